$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 105,28
$arr[0,0] = 0
$arr[0,1] = 6760228
$arr[0,2] = 'Germany Landesliga'
$arr[0,3] = 45088.41666666666
$arr[0,4] = 'SpVg Porz 1919'
$arr[0,5] = 'Bonn Endenich 1908'
$arr[0,6] = 4
$arr[0,7] = 1
$arr[0,8] = 'H'
$arr[0,9] = 2.25
$arr[0,10] = 3.5
$arr[0,11] = 2.625
$arr[0,12] = 1.75
$arr[0,13] = 3.6
$arr[0,14] = 3.8
$arr[0,15] = -0.5
$arr[0,16] = 1.8
$arr[0,17] = 2.0
$arr[0,18] = 2.75
$arr[0,19] = 1.8
$arr[0,20] = 2.0
$arr[0,21] = 0.75
$arr[0,22] = -1.0
$arr[0,23] = -1.0
$arr[0,24] = 0.8
$arr[0,25] = -1.0
$arr[0,26] = 0.8
$arr[0,27] = -1.0
$arr[1,0] = 1
$arr[1,1] = 6757276
$arr[1,2] = 'Germany Landesliga'
$arr[1,3] = 45088.41666666666
$arr[1,4] = 'BSC Rapid Chemnitz'
$arr[1,5] = 'FV Dresden 06 Laubegast'
$arr[1,6] = 3
$arr[1,7] = 0
$arr[1,8] = 'H'
$arr[1,9] = 2.6
$arr[1,10] = 3.5
$arr[1,11] = 2.25
$arr[1,12] = 2.6
$arr[1,13] = 3.5
$arr[1,14] = 2.25
$arr[1,15] = 0.25
$arr[1,16] = 1.75
$arr[1,17] = 2.05
$arr[1,18] = 3.25
$arr[1,19] = 1.775
$arr[1,20] = 2.025
$arr[1,21] = 1.6
$arr[1,22] = -1.0
$arr[1,23] = -1.0
$arr[1,24] = 0.75
$arr[1,25] = -1.0
$arr[1,26] = -0.5
$arr[1,27] = 0.5125
$arr[2,0] = 2
$arr[2,1] = 6781316
$arr[2,2] = 'Germany Landesliga'
$arr[2,3] = 45094.41666666666
$arr[2,4] = 'SV Schott Jena'
$arr[2,5] = 'SV 09 Arnstadt'
$arr[2,6] = 0
$arr[2,7] = 2
$arr[2,8] = 'A'
$arr[2,9] = 5.0
$arr[2,10] = 4.5
$arr[2,11] = 1.45
$arr[2,12] = 6.5
$arr[2,13] = 4.333
$arr[2,14] = 1.363
$arr[2,15] = 1.5
$arr[2,16] = 1.825
$arr[2,17] = 1.975
$arr[2,18] = 3.0
$arr[2,19] = 1.825
$arr[2,20] = 1.975
$arr[2,21] = -1.0
$arr[2,22] = -1.0
$arr[2,23] = 0.363
$arr[2,24] = -1.0
$arr[2,25] = 0.9750000000000001
$arr[2,26] = -1.0
$arr[2,27] = 0.9750000000000001
$arr[3,0] = 3
$arr[3,1] = 6781315
$arr[3,2] = 'Germany Landesliga'
$arr[3,3] = 45094.41666666666
$arr[3,4] = 'SSV Markranstadt'
$arr[3,5] = 'BSC Rapid Chemnitz'
$arr[3,6] = 2
$arr[3,7] = 0
$arr[3,8] = 'H'
$arr[3,9] = 1.25
$arr[3,10] = 4.75
$arr[3,11] = 10.0
$arr[3,12] = 1.222
$arr[3,13] = 5.25
$arr[3,14] = 8.5
$arr[3,15] = -2.0
$arr[3,16] = 1.925
$arr[3,17] = 1.875
$arr[3,18] = 3.5
$arr[3,19] = 1.775
$arr[3,20] = 1.925
$arr[3,21] = 0.222
$arr[3,22] = -1.0
$arr[3,23] = -1.0
$arr[3,24] = 0.0
$arr[3,25] = 0.0
$arr[3,26] = -1.0
$arr[3,27] = 0.925
$arr[4,0] = 4
$arr[4,1] = 7013106
$arr[4,2] = 'Germany Landesliga'
$arr[4,3] = 45142.60416666666
$arr[4,4] = 'ASV Suchteln'
$arr[4,5] = 'SC Victoria Mennrath'
$arr[4,6] = 3
$arr[4,7] = 4
$arr[4,8] = 'A'
$arr[4,9] = 1.833
$arr[4,10] = 4.333
$arr[4,11] = 3.0
$arr[4,12] = 1.615
$arr[4,13] = 4.75
$arr[4,14] = 3.5
$arr[4,15] = -1.0
$arr[4,16] = 1.975
$arr[4,17] = 1.825
$arr[4,18] = 4.25
$arr[4,19] = 1.975
$arr[4,20] = 1.825
$arr[4,21] = -1.0
$arr[4,22] = -1.0
$arr[4,23] = 2.5
$arr[4,24] = -1.0
$arr[4,25] = 0.825
$arr[4,26] = 0.9750000000000001
$arr[4,27] = -1.0
$arr[5,0] = 5
$arr[5,1] = 7018664
$arr[5,2] = 'Germany Landesliga'
$arr[5,3] = 45143.54166666666
$arr[5,4] = 'DJK Dilkrath'
$arr[5,5] = 'Holzheimer SG'
$arr[5,6] = 1
$arr[5,7] = 3
$arr[5,8] = 'A'
$arr[5,9] = 2.25
$arr[5,10] = 4.2
$arr[5,11] = 2.3
$arr[5,12] = 2.25
$arr[5,13] = 4.2
$arr[5,14] = 2.3
$arr[5,15] = 0.0
$arr[5,16] = 1.85
$arr[5,17] = 1.95
$arr[5,18] = 4.0
$arr[5,19] = 1.9
$arr[5,20] = 1.9
$arr[5,21] = -1.0
$arr[5,22] = -1.0
$arr[5,23] = 1.3
$arr[5,24] = -1.0
$arr[5,25] = 0.95
$arr[5,26] = 0.0
$arr[5,27] = 0.0
$arr[6,0] = 6
$arr[6,1] = 7021985
$arr[6,2] = 'Germany Landesliga'
$arr[6,3] = 45144.41666666666
$arr[6,4] = 'SC Velbert'
$arr[6,5] = 'SV Wermelskirchen'
$arr[6,6] = 2
$arr[6,7] = 0
$arr[6,8] = 'H'
$arr[6,9] = 1.833
$arr[6,10] = 4.5
$arr[6,11] = 2.875
$arr[6,12] = 1.7
$arr[6,13] = 4.75
$arr[6,14] = 3.2
$arr[6,15] = -0.75
$arr[6,16] = 1.95
$arr[6,17] = 1.85
$arr[6,18] = 3.25
$arr[6,19] = 1.8
$arr[6,20] = 2.0
$arr[6,21] = 0.7
$arr[6,22] = -1.0
$arr[6,23] = -1.0
$arr[6,24] = 0.95
$arr[6,25] = -1.0
$arr[6,26] = -1.0
$arr[6,27] = 1.0
$arr[7,0] = 7
$arr[7,1] = 7021984
$arr[7,2] = 'Germany Landesliga'
$arr[7,3] = 45144.4375
$arr[7,4] = 'TuRU Dsseldorf'
$arr[7,5] = 'SC KapellenErft'
$arr[7,6] = 0
$arr[7,7] = 1
$arr[7,8] = 'A'
$arr[7,9] = 1.666
$arr[7,10] = 4.5
$arr[7,11] = 3.4
$arr[7,12] = 1.666
$arr[7,13] = 4.5
$arr[7,14] = 3.4
$arr[7,15] = -0.75
$arr[7,16] = 1.875
$arr[7,17] = 1.925
$arr[7,18] = 3.5
$arr[7,19] = 1.9
$arr[7,20] = 1.9
$arr[7,21] = -1.0
$arr[7,22] = -1.0
$arr[7,23] = 2.4
$arr[7,24] = -1.0
$arr[7,25] = 0.925
$arr[7,26] = -1.0
$arr[7,27] = 0.8999999999999999
$arr[8,0] = 8
$arr[8,1] = 7035049
$arr[8,2] = 'Germany Landesliga'
$arr[8,3] = 45147.60416666666
$arr[8,4] = 'FC Kray'
$arr[8,5] = 'SV Honnepel Nierdermormter'
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 'H'
$arr[8,9] = 1.833
$arr[8,10] = 3.8
$arr[8,11] = 3.25
$arr[8,12] = 1.727
$arr[8,13] = 3.8
$arr[8,14] = 3.75
$arr[8,15] = -0.75
$arr[8,16] = 1.925
$arr[8,17] = 1.875
$arr[8,18] = 3.5
$arr[8,19] = 1.925
$arr[8,20] = 1.775
$arr[8,21] = 0.7270000000000001
$arr[8,22] = -1.0
$arr[8,23] = -1.0
$arr[8,24] = 0.925
$arr[8,25] = -1.0
$arr[8,26] = 0.925
$arr[8,27] = -1.0
$arr[9,0] = 9
$arr[9,1] = 7035048
$arr[9,2] = 'Germany Landesliga'
$arr[9,3] = 45147.625
$arr[9,4] = 'SG Unterrath'
$arr[9,5] = 'TuRU Dsseldorf'
$arr[9,6] = 1
$arr[9,7] = 0
$arr[9,8] = 'H'
$arr[9,9] = 3.25
$arr[9,10] = 4.0
$arr[9,11] = 1.8
$arr[9,12] = 2.9
$arr[9,13] = 4.0
$arr[9,14] = 1.95
$arr[9,15] = 0.5
$arr[9,16] = 1.8
$arr[9,17] = 2.0
$arr[9,18] = 3.0
$arr[9,19] = 1.75
$arr[9,20] = 1.95
$arr[9,21] = 1.9
$arr[9,22] = -1.0
$arr[9,23] = -1.0
$arr[9,24] = 0.8
$arr[9,25] = -1.0
$arr[9,26] = -1.0
$arr[9,27] = 0.95
$arr[10,0] = 10
$arr[10,1] = 7035047
$arr[10,2] = 'Germany Landesliga'
$arr[10,3] = 45147.625
$arr[10,4] = 'SC Dsseldorf West'
$arr[10,5] = 'VfL Viktoria JuchenGarzweiler'
$arr[10,6] = 3
$arr[10,7] = 4
$arr[10,8] = 'A'
$arr[10,9] = 1.909
$arr[10,10] = 3.75
$arr[10,11] = 3.1
$arr[10,12] = 2.2
$arr[10,13] = 3.6
$arr[10,14] = 2.625
$arr[10,15] = -0.25
$arr[10,16] = 2.0
$arr[10,17] = 1.8
$arr[10,18] = 3.0
$arr[10,19] = 1.825
$arr[10,20] = 1.975
$arr[10,21] = -1.0
$arr[10,22] = -1.0
$arr[10,23] = 1.625
$arr[10,24] = -1.0
$arr[10,25] = 0.8
$arr[10,26] = 0.825
$arr[10,27] = -1.0
$arr[11,0] = 11
$arr[11,1] = 7035046
$arr[11,2] = 'Germany Landesliga'
$arr[11,3] = 45147.625
$arr[11,4] = 'Cronenberger SC'
$arr[11,5] = 'FC Viersen'
$arr[11,6] = 0
$arr[11,7] = 2
$arr[11,8] = 'A'
$arr[11,9] = 2.0
$arr[11,10] = 3.6
$arr[11,11] = 3.0
$arr[11,12] = 2.0
$arr[11,13] = 3.6
$arr[11,14] = 3.0
$arr[11,15] = -0.25
$arr[11,16] = 1.8
$arr[11,17] = 2.0
$arr[11,18] = 2.75
$arr[11,19] = 1.8
$arr[11,20] = 2.0
$arr[11,21] = -1.0
$arr[11,22] = -1.0
$arr[11,23] = 2.0
$arr[11,24] = -1.0
$arr[11,25] = 1.0
$arr[11,26] = -1.0
$arr[11,27] = 1.0
$arr[12,0] = 12
$arr[12,1] = 7089910
$arr[12,2] = 'Germany Landesliga'
$arr[12,3] = 45158.41666666666
$arr[12,4] = 'ASV Mettmann'
$arr[12,5] = 'TuRU Dsseldorf'
$arr[12,6] = 2
$arr[12,7] = 1
$arr[12,8] = 'H'
$arr[12,9] = 3.25
$arr[12,10] = 4.0
$arr[12,11] = 1.8
$arr[12,12] = 3.25
$arr[12,13] = 4.0
$arr[12,14] = 1.8
$arr[12,15] = 0.5
$arr[12,16] = 1.975
$arr[12,17] = 1.825
$arr[12,18] = 3.25
$arr[12,19] = 1.85
$arr[12,20] = 1.95
$arr[12,21] = 2.25
$arr[12,22] = -1.0
$arr[12,23] = -1.0
$arr[12,24] = 0.9750000000000001
$arr[12,25] = -1.0
$arr[12,26] = -0.5
$arr[12,27] = 0.475
$arr[13,0] = 13
$arr[13,1] = 7089911
$arr[13,2] = 'Germany Landesliga'
$arr[13,3] = 45158.41666666666
$arr[13,4] = 'Spvgg Steele 0309'
$arr[13,5] = 'VfB Frohnhausen'
$arr[13,6] = 4
$arr[13,7] = 3
$arr[13,8] = 'H'
$arr[13,9] = 2.25
$arr[13,10] = 3.75
$arr[13,11] = 2.5
$arr[13,12] = 2.25
$arr[13,13] = 3.75
$arr[13,14] = 2.5
$arr[13,15] = 0.0
$arr[13,16] = 1.8
$arr[13,17] = 2.0
$arr[13,18] = 3.5
$arr[13,19] = 1.8
$arr[13,20] = 2.0
$arr[13,21] = 1.25
$arr[13,22] = -1.0
$arr[13,23] = -1.0
$arr[13,24] = 0.8
$arr[13,25] = -1.0
$arr[13,26] = 0.8
$arr[13,27] = -1.0
$arr[14,0] = 14
$arr[14,1] = 7117276
$arr[14,2] = 'Germany Landesliga'
$arr[14,3] = 45163.58333333334
$arr[14,4] = 'ESC Rellinghausen'
$arr[14,5] = 'DJK Arminia Klosterhardt'
$arr[14,6] = 2
$arr[14,7] = 0
$arr[14,8] = 'H'
$arr[14,9] = 2.3
$arr[14,10] = 3.6
$arr[14,11] = 2.5
$arr[14,12] = 2.3
$arr[14,13] = 3.6
$arr[14,14] = 2.5
$arr[14,15] = 0.0
$arr[14,16] = 1.8
$arr[14,17] = 2.0
$arr[14,18] = 3.0
$arr[14,19] = 1.825
$arr[14,20] = 1.975
$arr[14,21] = 1.3
$arr[14,22] = -1.0
$arr[14,23] = -1.0
$arr[14,24] = 0.8
$arr[14,25] = -1.0
$arr[14,26] = -1.0
$arr[14,27] = 0.9750000000000001
$arr[15,0] = 15
$arr[15,1] = 7142347
$arr[15,2] = 'Germany Landesliga'
$arr[15,3] = 45169.60416666666
$arr[15,4] = 'BlauWeiss Mintard'
$arr[15,5] = 'VfB Frohnhausen'
$arr[15,6] = 3
$arr[15,7] = 0
$arr[15,8] = 'H'
$arr[15,9] = 1.909
$arr[15,10] = 3.75
$arr[15,11] = 3.1
$arr[15,12] = 1.909
$arr[15,13] = 3.8
$arr[15,14] = 3.1
$arr[15,15] = -0.5
$arr[15,16] = 1.975
$arr[15,17] = 1.825
$arr[15,18] = 3.5
$arr[15,19] = 1.925
$arr[15,20] = 1.775
$arr[15,21] = 0.909
$arr[15,22] = -1.0
$arr[15,23] = -1.0
$arr[15,24] = 0.9750000000000001
$arr[15,25] = -1.0
$arr[15,26] = -1.0
$arr[15,27] = 0.7749999999999999
$arr[16,0] = 16
$arr[16,1] = 7157093
$arr[16,2] = 'Germany Landesliga'
$arr[16,3] = 45172.41666666666
$arr[16,4] = 'Spvgg Steele 0309'
$arr[16,5] = 'SV Honnepel Nierdermormter'
$arr[16,6] = 3
$arr[16,7] = 0
$arr[16,8] = 'H'
$arr[16,9] = 1.615
$arr[16,10] = 5.0
$arr[16,11] = 3.4
$arr[16,12] = 1.55
$arr[16,13] = 5.0
$arr[16,14] = 3.75
$arr[16,15] = -1.0
$arr[16,16] = 1.9
$arr[16,17] = 1.9
$arr[16,18] = 3.75
$arr[16,19] = 1.925
$arr[16,20] = 1.875
$arr[16,21] = 0.55
$arr[16,22] = -1.0
$arr[16,23] = -1.0
$arr[16,24] = 0.8999999999999999
$arr[16,25] = -1.0
$arr[16,26] = -1.0
$arr[16,27] = 0.875
$arr[17,0] = 17
$arr[17,1] = 7178367
$arr[17,2] = 'Germany Landesliga'
$arr[17,3] = 45177.625
$arr[17,4] = 'SV Scherpenberg'
$arr[17,5] = 'SGS EssenSchonebeck'
$arr[17,6] = 2
$arr[17,7] = 1
$arr[17,8] = 'H'
$arr[17,9] = 2.0
$arr[17,10] = 4.5
$arr[17,11] = 2.55
$arr[17,12] = 2.0
$arr[17,13] = 4.5
$arr[17,14] = 2.55
$arr[17,15] = -0.25
$arr[17,16] = 1.85
$arr[17,17] = 1.95
$arr[17,18] = 4.25
$arr[17,19] = 1.9
$arr[17,20] = 1.9
$arr[17,21] = 1.0
$arr[17,22] = -1.0
$arr[17,23] = -1.0
$arr[17,24] = 0.8500000000000001
$arr[17,25] = -1.0
$arr[17,26] = -1.0
$arr[17,27] = 0.8999999999999999
$arr[18,0] = 18
$arr[18,1] = 7189272
$arr[18,2] = 'Germany Landesliga'
$arr[18,3] = 45178.52083333334
$arr[18,4] = 'Holzheimer SG'
$arr[18,5] = 'SC KapellenErft'
$arr[18,6] = 0
$arr[18,7] = 3
$arr[18,8] = 'A'
$arr[18,9] = 2.25
$arr[18,10] = 4.0
$arr[18,11] = 2.4
$arr[18,12] = 2.05
$arr[18,13] = 4.2
$arr[18,14] = 2.625
$arr[18,15] = -0.25
$arr[18,16] = 1.875
$arr[18,17] = 1.925
$arr[18,18] = 3.5
$arr[18,19] = 1.95
$arr[18,20] = 1.85
$arr[18,21] = -1.0
$arr[18,22] = -1.0
$arr[18,23] = 1.625
$arr[18,24] = -1.0
$arr[18,25] = 0.925
$arr[18,26] = -1.0
$arr[18,27] = 0.8500000000000001
$arr[19,0] = 19
$arr[19,1] = 7212747
$arr[19,2] = 'Germany Landesliga'
$arr[19,3] = 45184.625
$arr[19,4] = 'SC Dsseldorf West'
$arr[19,5] = 'ASV Suchteln'
$arr[19,6] = 1
$arr[19,7] = 0
$arr[19,8] = 'H'
$arr[19,9] = 2.25
$arr[19,10] = 4.333
$arr[19,11] = 2.3
$arr[19,12] = 2.25
$arr[19,13] = 4.333
$arr[19,14] = 2.3
$arr[19,15] = 0.0
$arr[19,16] = 1.85
$arr[19,17] = 1.95
$arr[19,18] = 3.75
$arr[19,19] = 1.975
$arr[19,20] = 1.825
$arr[19,21] = 1.25
$arr[19,22] = -1.0
$arr[19,23] = -1.0
$arr[19,24] = 0.8500000000000001
$arr[19,25] = -1.0
$arr[19,26] = -1.0
$arr[19,27] = 0.825
$arr[20,0] = 20
$arr[20,1] = 7219907
$arr[20,2] = 'Germany Landesliga'
$arr[20,3] = 45185.54166666666
$arr[20,4] = 'Ahrweiler BC'
$arr[20,5] = 'Niederrossbach'
$arr[20,6] = 6
$arr[20,7] = 2
$arr[20,8] = 'H'
$arr[20,9] = 1.5
$arr[20,10] = 5.0
$arr[20,11] = 4.0
$arr[20,12] = 1.444
$arr[20,13] = 5.25
$arr[20,14] = 4.333
$arr[20,15] = -1.25
$arr[20,16] = 1.8
$arr[20,17] = 2.0
$arr[20,18] = 4.5
$arr[20,19] = 1.925
$arr[20,20] = 1.875
$arr[20,21] = 0.444
$arr[20,22] = -1.0
$arr[20,23] = -1.0
$arr[20,24] = 0.8
$arr[20,25] = -1.0
$arr[20,26] = 0.925
$arr[20,27] = -1.0
$arr[21,0] = 21
$arr[21,1] = 7221816
$arr[21,2] = 'Germany Landesliga'
$arr[21,3] = 45186.41666666666
$arr[21,4] = 'FC Monheim'
$arr[21,5] = 'SC Velbert'
$arr[21,6] = 1
$arr[21,7] = 1
$arr[21,8] = 'D'
$arr[21,9] = 1.8
$arr[21,10] = 3.4
$arr[21,11] = 3.8
$arr[21,12] = 1.8
$arr[21,13] = 3.4
$arr[21,14] = 3.8
$arr[21,15] = -0.5
$arr[21,16] = 1.85
$arr[21,17] = 1.95
$arr[21,18] = 2.5
$arr[21,19] = 2.0
$arr[21,20] = 1.8
$arr[21,21] = -1.0
$arr[21,22] = 2.4
$arr[21,23] = -1.0
$arr[21,24] = -1.0
$arr[21,25] = 0.95
$arr[21,26] = -1.0
$arr[21,27] = 0.8
$arr[22,0] = 22
$arr[22,1] = 7235837
$arr[22,2] = 'Germany Landesliga'
$arr[22,3] = 45190.625
$arr[22,4] = 'SV Helpenstein'
$arr[22,5] = 'Teutonia Weiden'
$arr[22,6] = 1
$arr[22,7] = 2
$arr[22,8] = 'A'
$arr[22,9] = 2.75
$arr[22,10] = 3.4
$arr[22,11] = 2.2
$arr[22,12] = 2.9
$arr[22,13] = 3.25
$arr[22,14] = 2.15
$arr[22,15] = 0.25
$arr[22,16] = 1.875
$arr[22,17] = 1.925
$arr[22,18] = 3.25
$arr[22,19] = 1.925
$arr[22,20] = 1.875
$arr[22,21] = -1.0
$arr[22,22] = -1.0
$arr[22,23] = 1.15
$arr[22,24] = -1.0
$arr[22,25] = 0.925
$arr[22,26] = -0.5
$arr[22,27] = 0.4375
$arr[23,0] = 23
$arr[23,1] = 7239956
$arr[23,2] = 'Germany Landesliga'
$arr[23,3] = 45191.60416666666
$arr[23,4] = 'VFB Hilden II'
$arr[23,5] = 'SC Dsseldorf West'
$arr[23,6] = 6
$arr[23,7] = 3
$arr[23,8] = 'H'
$arr[23,9] = 1.8
$arr[23,10] = 4.0
$arr[23,11] = 3.25
$arr[23,12] = 1.8
$arr[23,13] = 4.0
$arr[23,14] = 3.25
$arr[23,15] = -0.5
$arr[23,16] = 1.85
$arr[23,17] = 1.95
$arr[23,18] = 3.75
$arr[23,19] = 1.775
$arr[23,20] = 2.025
$arr[23,21] = 0.8
$arr[23,22] = -1.0
$arr[23,23] = -1.0
$arr[23,24] = 0.8500000000000001
$arr[23,25] = -1.0
$arr[23,26] = 0.7749999999999999
$arr[23,27] = -1.0
$arr[24,0] = 24
$arr[24,1] = 7244414
$arr[24,2] = 'Germany Landesliga'
$arr[24,3] = 45192.60416666666
$arr[24,4] = 'SC Victoria Mennrath'
$arr[24,5] = 'SC KapellenErft'
$arr[24,6] = 0
$arr[24,7] = 3
$arr[24,8] = 'A'
$arr[24,9] = 5.0
$arr[24,10] = 3.75
$arr[24,11] = 1.55
$arr[24,12] = 5.5
$arr[24,13] = 4.2
$arr[24,14] = 1.444
$arr[24,15] = 1.25
$arr[24,16] = 1.975
$arr[24,17] = 1.825
$arr[24,18] = 4.0
$arr[24,19] = 2.0
$arr[24,20] = 1.8
$arr[24,21] = -1.0
$arr[24,22] = -1.0
$arr[24,23] = 0.444
$arr[24,24] = -1.0
$arr[24,25] = 0.825
$arr[24,26] = -1.0
$arr[24,27] = 0.8
$arr[25,0] = 25
$arr[25,1] = 7271358
$arr[25,2] = 'Germany Landesliga'
$arr[25,3] = 45198.58333333334
$arr[25,4] = 'Spvgg Steele 0309'
$arr[25,5] = 'FC Kray'
$arr[25,6] = 2
$arr[25,7] = 3
$arr[25,8] = 'A'
$arr[25,9] = 3.4
$arr[25,10] = 4.333
$arr[25,11] = 1.7
$arr[25,12] = 3.3
$arr[25,13] = 4.333
$arr[25,14] = 1.7
$arr[25,15] = 0.75
$arr[25,16] = 1.875
$arr[25,17] = 1.925
$arr[25,18] = 4.0
$arr[25,19] = 1.85
$arr[25,20] = 1.95
$arr[25,21] = -1.0
$arr[25,22] = -1.0
$arr[25,23] = 0.7
$arr[25,24] = -0.5
$arr[25,25] = 0.4625
$arr[25,26] = 0.8500000000000001
$arr[25,27] = -1.0
$arr[26,0] = 26
$arr[26,1] = 7280312
$arr[26,2] = 'Germany Landesliga'
$arr[26,3] = 45200.41666666666
$arr[26,4] = 'FC Monheim'
$arr[26,5] = 'ASV Suchteln'
$arr[26,6] = 1
$arr[26,7] = 1
$arr[26,8] = 'D'
$arr[26,9] = 1.4
$arr[26,10] = 4.333
$arr[26,11] = 6.0
$arr[26,12] = 1.4
$arr[26,13] = 4.333
$arr[26,14] = 6.0
$arr[26,15] = -1.25
$arr[26,16] = 1.8
$arr[26,17] = 2.0
$arr[26,18] = 3.75
$arr[26,19] = 2.0
$arr[26,20] = 1.8
$arr[26,21] = -1.0
$arr[26,22] = 3.333
$arr[26,23] = -1.0
$arr[26,24] = -1.0
$arr[26,25] = 1.0
$arr[26,26] = -1.0
$arr[26,27] = 0.8
$arr[27,0] = 27
$arr[27,1] = 7280311
$arr[27,2] = 'Germany Landesliga'
$arr[27,3] = 45200.42708333334
$arr[27,4] = 'SV Schlebusch'
$arr[27,5] = 'TUS 05 Oberpleis'
$arr[27,6] = 2
$arr[27,7] = 4
$arr[27,8] = 'A'
$arr[27,9] = 2.7
$arr[27,10] = 3.75
$arr[27,11] = 2.1
$arr[27,12] = 2.7
$arr[27,13] = 3.75
$arr[27,14] = 2.1
$arr[27,15] = 0.25
$arr[27,16] = 1.875
$arr[27,17] = 1.925
$arr[27,18] = 3.25
$arr[27,19] = 2.0
$arr[27,20] = 1.8
$arr[27,21] = -1.0
$arr[27,22] = -1.0
$arr[27,23] = 1.1
$arr[27,24] = -1.0
$arr[27,25] = 0.925
$arr[27,26] = 1.0
$arr[27,27] = -1.0
$arr[28,0] = 28
$arr[28,1] = 7325273
$arr[28,2] = 'Germany Landesliga'
$arr[28,3] = 45210.625
$arr[28,4] = 'SSV Homburg Numbrecht'
$arr[28,5] = 'SC Fortuna Bonn'
$arr[28,6] = 1
$arr[28,7] = 1
$arr[28,8] = 'D'
$arr[28,9] = 2.0
$arr[28,10] = 3.75
$arr[28,11] = 2.9
$arr[28,12] = 2.0
$arr[28,13] = 3.75
$arr[28,14] = 2.9
$arr[28,15] = -0.25
$arr[28,16] = 1.825
$arr[28,17] = 1.975
$arr[28,18] = 3.75
$arr[28,19] = 1.975
$arr[28,20] = 1.825
$arr[28,21] = -1.0
$arr[28,22] = 2.75
$arr[28,23] = -1.0
$arr[28,24] = -0.5
$arr[28,25] = 0.4875
$arr[28,26] = -1.0
$arr[28,27] = 0.825
$arr[29,0] = 29
$arr[29,1] = 7338519
$arr[29,2] = 'Germany Landesliga'
$arr[29,3] = 45213.375
$arr[29,4] = 'FSV 06 Ohratal'
$arr[29,5] = 'SpVgg Geratal'
$arr[29,6] = 0
$arr[29,7] = 0
$arr[29,8] = 'D'
$arr[29,9] = 2.2
$arr[29,10] = 4.2
$arr[29,11] = 2.375
$arr[29,12] = 2.2
$arr[29,13] = 4.2
$arr[29,14] = 2.375
$arr[29,15] = 0.0
$arr[29,16] = 1.8
$arr[29,17] = 2.0
$arr[29,18] = 3.75
$arr[29,19] = 1.925
$arr[29,20] = 1.875
$arr[29,21] = -1.0
$arr[29,22] = 3.2
$arr[29,23] = -1.0
$arr[29,24] = 0.0
$arr[29,25] = 0.0
$arr[29,26] = -1.0
$arr[29,27] = 0.875
$arr[30,0] = 30
$arr[30,1] = 7343910
$arr[30,2] = 'Germany Landesliga'
$arr[30,3] = 45214.41666666666
$arr[30,4] = 'FC Monheim'
$arr[30,5] = 'DJK Dilkrath'
$arr[30,6] = 5
$arr[30,7] = 0
$arr[30,8] = 'H'
$arr[30,9] = 1.222
$arr[30,10] = 6.0
$arr[30,11] = 7.5
$arr[30,12] = 1.222
$arr[30,13] = 6.0
$arr[30,14] = 7.0
$arr[30,15] = -1.75
$arr[30,16] = 1.75
$arr[30,17] = 2.05
$arr[30,18] = 3.75
$arr[30,19] = 2.0
$arr[30,20] = 1.8
$arr[30,21] = 0.222
$arr[30,22] = -1.0
$arr[30,23] = -1.0
$arr[30,24] = 0.75
$arr[30,25] = -1.0
$arr[30,26] = 1.0
$arr[30,27] = -1.0
$arr[31,0] = 31
$arr[31,1] = 7354380
$arr[31,2] = 'Germany Landesliga'
$arr[31,3] = 45217.625
$arr[31,4] = 'FC Pesch 1956'
$arr[31,5] = 'FSV NeunkirchenSeelscheid'
$arr[31,6] = 2
$arr[31,7] = 1
$arr[31,8] = 'H'
$arr[31,9] = 1.5
$arr[31,10] = 4.0
$arr[31,11] = 5.0
$arr[31,12] = 1.4
$arr[31,13] = 4.333
$arr[31,14] = 5.5
$arr[31,15] = -1.5
$arr[31,16] = 1.975
$arr[31,17] = 1.825
$arr[31,18] = 3.75
$arr[31,19] = 1.975
$arr[31,20] = 1.825
$arr[31,21] = 0.3999999999999999
$arr[31,22] = -1.0
$arr[31,23] = -1.0
$arr[31,24] = -1.0
$arr[31,25] = 0.825
$arr[31,26] = -1.0
$arr[31,27] = 0.825
$arr[32,0] = 32
$arr[32,1] = 7357613
$arr[32,2] = 'Germany Landesliga'
$arr[32,3] = 45218.60416666666
$arr[32,4] = 'SV Eilendorf'
$arr[32,5] = 'Teutonia Weiden'
$arr[32,6] = 2
$arr[32,7] = 3
$arr[32,8] = 'A'
$arr[32,9] = 5.0
$arr[32,10] = 4.0
$arr[32,11] = 1.5
$arr[32,12] = 6.0
$arr[32,13] = 4.75
$arr[32,14] = 1.333
$arr[32,15] = 1.75
$arr[32,16] = 1.9
$arr[32,17] = 1.9
$arr[32,18] = 4.25
$arr[32,19] = 1.925
$arr[32,20] = 1.875
$arr[32,21] = -1.0
$arr[32,22] = -1.0
$arr[32,23] = 0.333
$arr[32,24] = 0.8999999999999999
$arr[32,25] = -1.0
$arr[32,26] = 0.925
$arr[32,27] = -1.0
$arr[33,0] = 33
$arr[33,1] = 7371459
$arr[33,2] = 'Germany Landesliga'
$arr[33,3] = 45221.4375
$arr[33,4] = 'VFB Hilden II'
$arr[33,5] = 'TuRU Dsseldorf'
$arr[33,6] = 0
$arr[33,7] = 1
$arr[33,8] = 'A'
$arr[33,9] = 2.0
$arr[33,10] = 4.2
$arr[33,11] = 2.7
$arr[33,12] = 2.2
$arr[33,13] = 4.2
$arr[33,14] = 2.4
$arr[33,15] = 0.0
$arr[33,16] = 1.8
$arr[33,17] = 2.0
$arr[33,18] = 3.25
$arr[33,19] = 1.8
$arr[33,20] = 2.0
$arr[33,21] = -1.0
$arr[33,22] = -1.0
$arr[33,23] = 1.4
$arr[33,24] = -1.0
$arr[33,25] = 1.0
$arr[33,26] = -1.0
$arr[33,27] = 1.0
$arr[34,0] = 34
$arr[34,1] = 7388066
$arr[34,2] = 'Germany Landesliga'
$arr[34,3] = 45225.60416666666
$arr[34,4] = 'SV Eintracht Verlautenheide'
$arr[34,5] = 'Arminia Eilendorf'
$arr[34,6] = 4
$arr[34,7] = 2
$arr[34,8] = 'H'
$arr[34,9] = 1.666
$arr[34,10] = 4.0
$arr[34,11] = 3.75
$arr[34,12] = 1.666
$arr[34,13] = 4.0
$arr[34,14] = 3.75
$arr[34,15] = -0.75
$arr[34,16] = 1.875
$arr[34,17] = 1.925
$arr[34,18] = 3.5
$arr[34,19] = 1.975
$arr[34,20] = 1.825
$arr[34,21] = 0.6659999999999999
$arr[34,22] = -1.0
$arr[34,23] = -1.0
$arr[34,24] = 0.875
$arr[34,25] = -1.0
$arr[34,26] = 0.9750000000000001
$arr[34,27] = -1.0
$arr[35,0] = 35
$arr[35,1] = 7398238
$arr[35,2] = 'Germany Landesliga'
$arr[35,3] = 45228.45833333334
$arr[35,4] = 'FC Monheim'
$arr[35,5] = 'VSF Amern'
$arr[35,6] = 2
$arr[35,7] = 3
$arr[35,8] = 'A'
$arr[35,9] = 1.333
$arr[35,10] = 5.0
$arr[35,11] = 6.0
$arr[35,12] = 1.285
$arr[35,13] = 5.25
$arr[35,14] = 7.0
$arr[35,15] = -1.75
$arr[35,16] = 1.8
$arr[35,17] = 2.0
$arr[35,18] = 4.0
$arr[35,19] = 1.9
$arr[35,20] = 1.9
$arr[35,21] = -1.0
$arr[35,22] = -1.0
$arr[35,23] = 6.0
$arr[35,24] = -1.0
$arr[35,25] = 1.0
$arr[35,26] = 0.8999999999999999
$arr[35,27] = -1.0
$arr[36,0] = 36
$arr[36,1] = 7398237
$arr[36,2] = 'Germany Landesliga'
$arr[36,3] = 45228.46875
$arr[36,4] = 'SV Schlebusch'
$arr[36,5] = 'SSV Bornheim'
$arr[36,6] = 1
$arr[36,7] = 1
$arr[36,8] = 'D'
$arr[36,9] = 4.2
$arr[36,10] = 4.2
$arr[36,11] = 1.571
$arr[36,12] = 4.2
$arr[36,13] = 4.2
$arr[36,14] = 1.571
$arr[36,15] = 1.0
$arr[36,16] = 1.9
$arr[36,17] = 1.9
$arr[36,18] = 3.75
$arr[36,19] = 1.925
$arr[36,20] = 1.875
$arr[36,21] = -1.0
$arr[36,22] = 3.2
$arr[36,23] = -1.0
$arr[36,24] = 0.8999999999999999
$arr[36,25] = -1.0
$arr[36,26] = -1.0
$arr[36,27] = 0.875
$arr[37,0] = 37
$arr[37,1] = 7426706
$arr[37,2] = 'Germany Landesliga'
$arr[37,3] = 45234.375
$arr[37,4] = 'SV Schott Jena'
$arr[37,5] = 'FSV Schleiz'
$arr[37,6] = 0
$arr[37,7] = 0
$arr[37,8] = 'D'
$arr[37,9] = 2.4
$arr[37,10] = 3.5
$arr[37,11] = 2.4
$arr[37,12] = 2.4
$arr[37,13] = 3.5
$arr[37,14] = 2.4
$arr[37,15] = 0.0
$arr[37,16] = 1.9
$arr[37,17] = 1.9
$arr[37,18] = 4.0
$arr[37,19] = 1.9
$arr[37,20] = 1.9
$arr[37,21] = -1.0
$arr[37,22] = 2.5
$arr[37,23] = -1.0
$arr[37,24] = 0.0
$arr[37,25] = 0.0
$arr[37,26] = -1.0
$arr[37,27] = 0.8999999999999999
$arr[38,0] = 38
$arr[38,1] = 7465503
$arr[38,2] = 'Germany Landesliga'
$arr[38,3] = 45242.45833333334
$arr[38,4] = 'FC Monheim'
$arr[38,5] = 'FC Remscheid'
$arr[38,6] = 2
$arr[38,7] = 0
$arr[38,8] = 'H'
$arr[38,9] = 1.909
$arr[38,10] = 3.75
$arr[38,11] = 3.1
$arr[38,12] = 1.615
$arr[38,13] = 4.0
$arr[38,14] = 4.2
$arr[38,15] = -0.75
$arr[38,16] = 1.825
$arr[38,17] = 1.975
$arr[38,18] = 3.25
$arr[38,19] = 1.95
$arr[38,20] = 1.85
$arr[38,21] = 0.615
$arr[38,22] = -1.0
$arr[38,23] = -1.0
$arr[38,24] = 0.825
$arr[38,25] = -1.0
$arr[38,26] = -1.0
$arr[38,27] = 0.8500000000000001
$arr[39,0] = 39
$arr[39,1] = 7465502
$arr[39,2] = 'Germany Landesliga'
$arr[39,3] = 45242.45833333334
$arr[39,4] = 'FC Pesch 1956'
$arr[39,5] = 'SSV Bornheim'
$arr[39,6] = 4
$arr[39,7] = 0
$arr[39,8] = 'H'
$arr[39,9] = 2.0
$arr[39,10] = 3.75
$arr[39,11] = 2.875
$arr[39,12] = 2.0
$arr[39,13] = 3.75
$arr[39,14] = 2.9
$arr[39,15] = -0.25
$arr[39,16] = 1.8
$arr[39,17] = 2.0
$arr[39,18] = 3.0
$arr[39,19] = 1.775
$arr[39,20] = 2.025
$arr[39,21] = 1.0
$arr[39,22] = -1.0
$arr[39,23] = -1.0
$arr[39,24] = 0.8
$arr[39,25] = -1.0
$arr[39,26] = 0.7749999999999999
$arr[39,27] = -1.0
$arr[40,0] = 40
$arr[40,1] = 7485798
$arr[40,2] = 'Germany Landesliga'
$arr[40,3] = 45247.66666666666
$arr[40,4] = 'SC Victoria Mennrath'
$arr[40,5] = 'TuRU Dsseldorf'
$arr[40,6] = 4
$arr[40,7] = 0
$arr[40,8] = 'H'
$arr[40,9] = 2.7
$arr[40,10] = 4.0
$arr[40,11] = 2.025
$arr[40,12] = 3.6
$arr[40,13] = 4.333
$arr[40,14] = 1.65
$arr[40,15] = 0.75
$arr[40,16] = 1.975
$arr[40,17] = 1.825
$arr[40,18] = 3.75
$arr[40,19] = 1.925
$arr[40,20] = 1.875
$arr[40,21] = 2.6
$arr[40,22] = -1.0
$arr[40,23] = -1.0
$arr[40,24] = 0.9750000000000001
$arr[40,25] = -1.0
$arr[40,26] = 0.4625
$arr[40,27] = -0.5
$arr[41,0] = 41
$arr[41,1] = 7492563
$arr[41,2] = 'Germany Landesliga'
$arr[41,3] = 45249.47916666666
$arr[41,4] = 'MSV Dsseldorf'
$arr[41,5] = 'VFB Hilden II'
$arr[41,6] = 0
$arr[41,7] = 2
$arr[41,8] = 'A'
$arr[41,9] = 4.5
$arr[41,10] = 4.5
$arr[41,11] = 1.5
$arr[41,12] = 4.5
$arr[41,13] = 4.5
$arr[41,14] = 1.5
$arr[41,15] = 1.0
$arr[41,16] = 2.0
$arr[41,17] = 1.8
$arr[41,18] = 3.25
$arr[41,19] = 1.9
$arr[41,20] = 1.9
$arr[41,21] = -1.0
$arr[41,22] = -1.0
$arr[41,23] = 0.5
$arr[41,24] = -1.0
$arr[41,25] = 0.8
$arr[41,26] = -1.0
$arr[41,27] = 0.8999999999999999
$arr[42,0] = 42
$arr[42,1] = 7505166
$arr[42,2] = 'Germany Landesliga'
$arr[42,3] = 45252.66666666666
$arr[42,4] = 'SV HolthausenBiene'
$arr[42,5] = 'SV Bevern'
$arr[42,6] = 6
$arr[42,7] = 1
$arr[42,8] = 'H'
$arr[42,9] = 1.8
$arr[42,10] = 3.8
$arr[42,11] = 3.4
$arr[42,12] = 1.45
$arr[42,13] = 4.333
$arr[42,14] = 5.25
$arr[42,15] = -1.25
$arr[42,16] = 1.95
$arr[42,17] = 1.85
$arr[42,18] = 3.5
$arr[42,19] = 1.95
$arr[42,20] = 1.85
$arr[42,21] = 0.45
$arr[42,22] = -1.0
$arr[42,23] = -1.0
$arr[42,24] = 0.95
$arr[42,25] = -1.0
$arr[42,26] = 0.95
$arr[42,27] = -1.0
$arr[43,0] = 43
$arr[43,1] = 7511939
$arr[43,2] = 'Germany Landesliga'
$arr[43,3] = 45254.64583333334
$arr[43,4] = 'FC Monheim'
$arr[43,5] = 'SC KapellenErft'
$arr[43,6] = 3
$arr[43,7] = 1
$arr[43,8] = 'H'
$arr[43,9] = 2.2
$arr[43,10] = 3.5
$arr[43,11] = 2.7
$arr[43,12] = 1.8
$arr[43,13] = 3.5
$arr[43,14] = 3.6
$arr[43,15] = -0.5
$arr[43,16] = 1.85
$arr[43,17] = 1.95
$arr[43,18] = 3.25
$arr[43,19] = 1.85
$arr[43,20] = 1.95
$arr[43,21] = 0.8
$arr[43,22] = -1.0
$arr[43,23] = -1.0
$arr[43,24] = 0.8500000000000001
$arr[43,25] = -1.0
$arr[43,26] = 0.8500000000000001
$arr[43,27] = -1.0
$arr[44,0] = 44
$arr[44,1] = 7511941
$arr[44,2] = 'Germany Landesliga'
$arr[44,3] = 45254.66666666666
$arr[44,4] = 'SpVgg SterkradeNord'
$arr[44,5] = 'BlauWeiss Mintard'
$arr[44,6] = 3
$arr[44,7] = 1
$arr[44,8] = 'H'
$arr[44,9] = 3.4
$arr[44,10] = 4.2
$arr[44,11] = 1.727
$arr[44,12] = 2.8
$arr[44,13] = 4.2
$arr[44,14] = 1.95
$arr[44,15] = 0.5
$arr[44,16] = 1.8
$arr[44,17] = 2.0
$arr[44,18] = 3.5
$arr[44,19] = 1.875
$arr[44,20] = 1.925
$arr[44,21] = 1.8
$arr[44,22] = -1.0
$arr[44,23] = -1.0
$arr[44,24] = 0.8
$arr[44,25] = -1.0
$arr[44,26] = 0.875
$arr[44,27] = -1.0
$arr[45,0] = 45
$arr[45,1] = 7511940
$arr[45,2] = 'Germany Landesliga'
$arr[45,3] = 45254.66666666666
$arr[45,4] = 'ASV Suchteln'
$arr[45,5] = 'Holzheimer SG'
$arr[45,6] = 3
$arr[45,7] = 3
$arr[45,8] = 'D'
$arr[45,9] = 2.75
$arr[45,10] = 3.6
$arr[45,11] = 2.1
$arr[45,12] = 4.0
$arr[45,13] = 4.0
$arr[45,14] = 1.6
$arr[45,15] = 1.0
$arr[45,16] = 1.8
$arr[45,17] = 2.0
$arr[45,18] = 3.5
$arr[45,19] = 1.85
$arr[45,20] = 1.95
$arr[45,21] = -1.0
$arr[45,22] = 3.0
$arr[45,23] = -1.0
$arr[45,24] = 0.8
$arr[45,25] = -1.0
$arr[45,26] = 0.8500000000000001
$arr[45,27] = -1.0
$arr[46,0] = 46
$arr[46,1] = 7521305
$arr[46,2] = 'Germany Landesliga'
$arr[46,3] = 45256.58333333334
$arr[46,4] = 'DJK Rasensport AachenBrand'
$arr[46,5] = 'Kohlscheider BC'
$arr[46,6] = 2
$arr[46,7] = 1
$arr[46,8] = 'H'
$arr[46,9] = 2.0
$arr[46,10] = 4.5
$arr[46,11] = 2.55
$arr[46,12] = 1.909
$arr[46,13] = 4.5
$arr[46,14] = 2.7
$arr[46,15] = -0.5
$arr[46,16] = 1.975
$arr[46,17] = 1.825
$arr[46,18] = 4.0
$arr[46,19] = 1.85
$arr[46,20] = 1.95
$arr[46,21] = 0.909
$arr[46,22] = -1.0
$arr[46,23] = -1.0
$arr[46,24] = 0.9750000000000001
$arr[46,25] = -1.0
$arr[46,26] = -1.0
$arr[46,27] = 0.95
$arr[47,0] = 47
$arr[47,1] = 7538685
$arr[47,2] = 'Germany Landesliga'
$arr[47,3] = 45261.625
$arr[47,4] = 'Spvgg Steele 0309'
$arr[47,5] = 'SV Budberg'
$arr[47,6] = 1
$arr[47,7] = 3
$arr[47,8] = 'A'
$arr[47,9] = 2.375
$arr[47,10] = 4.2
$arr[47,11] = 2.2
$arr[47,12] = 2.4
$arr[47,13] = 4.333
$arr[47,14] = 2.15
$arr[47,15] = 0.25
$arr[47,16] = 1.8
$arr[47,17] = 2.0
$arr[47,18] = 3.75
$arr[47,19] = 1.95
$arr[47,20] = 1.85
$arr[47,21] = -1.0
$arr[47,22] = -1.0
$arr[47,23] = 1.15
$arr[47,24] = -1.0
$arr[47,25] = 1.0
$arr[47,26] = 0.475
$arr[47,27] = -0.5
$arr[48,0] = 48
$arr[48,1] = 7559508
$arr[48,2] = 'Germany Landesliga'
$arr[48,3] = 45266.64583333334
$arr[48,4] = 'FC Kray'
$arr[48,5] = 'SF Niederwenigern'
$arr[48,6] = 2
$arr[48,7] = 3
$arr[48,8] = 'A'
$arr[48,9] = 2.1
$arr[48,10] = 3.6
$arr[48,11] = 2.8
$arr[48,12] = 2.1
$arr[48,13] = 3.6
$arr[48,14] = 2.8
$arr[48,15] = -0.25
$arr[48,16] = 1.9
$arr[48,17] = 1.9
$arr[48,18] = 3.5
$arr[48,19] = 1.875
$arr[48,20] = 1.925
$arr[48,21] = -1.0
$arr[48,22] = -1.0
$arr[48,23] = 1.8
$arr[48,24] = -1.0
$arr[48,25] = 0.8999999999999999
$arr[48,26] = 0.875
$arr[48,27] = -1.0
$arr[49,0] = 49
$arr[49,1] = 7562219
$arr[49,2] = 'Germany Landesliga'
$arr[49,3] = 45267.67013888889
$arr[49,4] = 'FV Wiehl 2000'
$arr[49,5] = 'SC Rheinbach'
$arr[49,6] = 1
$arr[49,7] = 2
$arr[49,8] = 'A'
$arr[49,9] = 2.1
$arr[49,10] = 4.0
$arr[49,11] = 2.6
$arr[49,12] = 2.0
$arr[49,13] = 3.8
$arr[49,14] = 2.8
$arr[49,15] = -0.25
$arr[49,16] = 1.85
$arr[49,17] = 1.95
$arr[49,18] = 3.25
$arr[49,19] = 2.0
$arr[49,20] = 1.8
$arr[49,21] = -1.0
$arr[49,22] = -1.0
$arr[49,23] = 1.8
$arr[49,24] = -1.0
$arr[49,25] = 0.95
$arr[49,26] = -0.5
$arr[49,27] = 0.4
$arr[50,0] = 50
$arr[50,1] = 7566096
$arr[50,2] = 'Germany Landesliga'
$arr[50,3] = 45268.66666666666
$arr[50,4] = 'SC Victoria Mennrath'
$arr[50,5] = 'ASV Suchteln'
$arr[50,6] = 5
$arr[50,7] = 2
$arr[50,8] = 'H'
$arr[50,9] = 1.8
$arr[50,10] = 4.5
$arr[50,11] = 3.0
$arr[50,12] = 1.8
$arr[50,13] = 4.5
$arr[50,14] = 2.9
$arr[50,15] = -0.5
$arr[50,16] = 1.85
$arr[50,17] = 1.95
$arr[50,18] = 3.75
$arr[50,19] = 1.8
$arr[50,20] = 2.0
$arr[50,21] = 0.8
$arr[50,22] = -1.0
$arr[50,23] = -1.0
$arr[50,24] = 0.8500000000000001
$arr[50,25] = -1.0
$arr[50,26] = 0.8
$arr[50,27] = -1.0
$arr[51,0] = 51
$arr[51,1] = 7583056
$arr[51,2] = 'Germany Landesliga'
$arr[51,3] = 45273.66666666666
$arr[51,4] = 'DJK Arminia Klosterhardt'
$arr[51,5] = 'SF Niederwenigern'
$arr[51,6] = 2
$arr[51,7] = 2
$arr[51,8] = 'D'
$arr[51,9] = 3.6
$arr[51,10] = 4.2
$arr[51,11] = 1.666
$arr[51,12] = 3.8
$arr[51,13] = 4.5
$arr[51,14] = 1.571
$arr[51,15] = 0.75
$arr[51,16] = 2.075
$arr[51,17] = 1.725
$arr[51,18] = 3.75
$arr[51,19] = 1.9
$arr[51,20] = 1.9
$arr[51,21] = -1.0
$arr[51,22] = 3.5
$arr[51,23] = -1.0
$arr[51,24] = 1.075
$arr[51,25] = -1.0
$arr[51,26] = 0.45
$arr[51,27] = -0.5
$arr[52,0] = 52
$arr[52,1] = 7598132
$arr[52,2] = 'Germany Landesliga'
$arr[52,3] = 45277.47916666666
$arr[52,4] = 'MSV Dsseldorf'
$arr[52,5] = 'Holzheimer SG'
$arr[52,6] = 3
$arr[52,7] = 1
$arr[52,8] = 'H'
$arr[52,9] = 4.2
$arr[52,10] = 4.5
$arr[52,11] = 1.533
$arr[52,12] = 4.2
$arr[52,13] = 4.5
$arr[52,14] = 1.533
$arr[52,15] = 1.0
$arr[52,16] = 1.975
$arr[52,17] = 1.825
$arr[52,18] = 3.5
$arr[52,19] = 1.825
$arr[52,20] = 1.975
$arr[52,21] = 3.2
$arr[52,22] = -1.0
$arr[52,23] = -1.0
$arr[52,24] = 0.9750000000000001
$arr[52,25] = -1.0
$arr[52,26] = 0.825
$arr[52,27] = -1.0
$arr[53,0] = 53
$arr[53,1] = 7783594
$arr[53,2] = 'Germany Landesliga'
$arr[53,3] = 45326.4375
$arr[53,4] = 'Cronenberger SC'
$arr[53,5] = 'VfL Viktoria JuchenGarzweiler'
$arr[53,6] = 0
$arr[53,7] = 2
$arr[53,8] = 'A'
$arr[53,9] = 2.5
$arr[53,10] = 3.75
$arr[53,11] = 2.25
$arr[53,12] = 2.875
$arr[53,13] = 3.75
$arr[53,14] = 2.0
$arr[53,15] = 0.25
$arr[53,16] = 1.975
$arr[53,17] = 1.825
$arr[53,18] = 3.5
$arr[53,19] = 1.85
$arr[53,20] = 1.95
$arr[53,21] = -1.0
$arr[53,22] = -1.0
$arr[53,23] = 1.0
$arr[53,24] = -1.0
$arr[53,25] = 0.825
$arr[53,26] = -1.0
$arr[53,27] = 0.95
$arr[54,0] = 54
$arr[54,1] = 7847639
$arr[54,2] = 'Germany Landesliga'
$arr[54,3] = 45340.41666666666
$arr[54,4] = 'VfB Fortuna Chemnitz'
$arr[54,5] = 'SG Taucha 99'
$arr[54,6] = 0
$arr[54,7] = 2
$arr[54,8] = 'A'
$arr[54,9] = 2.15
$arr[54,10] = 3.75
$arr[54,11] = 2.6
$arr[54,12] = 2.15
$arr[54,13] = 3.75
$arr[54,14] = 2.625
$arr[54,15] = -0.25
$arr[54,16] = 1.975
$arr[54,17] = 1.825
$arr[54,18] = 2.75
$arr[54,19] = 1.95
$arr[54,20] = 1.85
$arr[54,21] = -1.0
$arr[54,22] = -1.0
$arr[54,23] = 1.625
$arr[54,24] = -1.0
$arr[54,25] = 0.825
$arr[54,26] = -1.0
$arr[54,27] = 0.8500000000000001
$arr[55,0] = 55
$arr[55,1] = 7865257
$arr[55,2] = 'Germany Landesliga'
$arr[55,3] = 45344.64583333334
$arr[55,4] = 'Arminia Eilendorf'
$arr[55,5] = 'Kohlscheider BC'
$arr[55,6] = 3
$arr[55,7] = 0
$arr[55,8] = 'H'
$arr[55,9] = 3.5
$arr[55,10] = 4.0
$arr[55,11] = 1.727
$arr[55,12] = 2.7
$arr[55,13] = 4.2
$arr[55,14] = 2.0
$arr[55,15] = 0.25
$arr[55,16] = 1.95
$arr[55,17] = 1.85
$arr[55,18] = 4.25
$arr[55,19] = 1.85
$arr[55,20] = 1.95
$arr[55,21] = 1.7
$arr[55,22] = -1.0
$arr[55,23] = -1.0
$arr[55,24] = 0.95
$arr[55,25] = -1.0
$arr[55,26] = -1.0
$arr[55,27] = 0.95
$arr[56,0] = 56
$arr[56,1] = 7871297
$arr[56,2] = 'Germany Landesliga'
$arr[56,3] = 45345.64583333334
$arr[56,4] = 'TuRU Dsseldorf'
$arr[56,5] = 'ASV Mettmann'
$arr[56,6] = 2
$arr[56,7] = 0
$arr[56,8] = 'H'
$arr[56,9] = 1.6
$arr[56,10] = 4.0
$arr[56,11] = 4.2
$arr[56,12] = 1.6
$arr[56,13] = 4.0
$arr[56,14] = 4.2
$arr[56,15] = -1.0
$arr[56,16] = 2.025
$arr[56,17] = 1.775
$arr[56,18] = 3.0
$arr[56,19] = 1.925
$arr[56,20] = 1.875
$arr[56,21] = 0.6000000000000001
$arr[56,22] = -1.0
$arr[56,23] = -1.0
$arr[56,24] = 1.025
$arr[56,25] = -1.0
$arr[56,26] = -1.0
$arr[56,27] = 0.875
$arr[57,0] = 57
$arr[57,1] = 7873626
$arr[57,2] = 'Germany Landesliga'
$arr[57,3] = 45346.41666666666
$arr[57,4] = 'SG Taucha 99'
$arr[57,5] = 'VfL PirnaCopitz'
$arr[57,6] = 3
$arr[57,7] = 0
$arr[57,8] = 'H'
$arr[57,9] = 1.333
$arr[57,10] = 5.0
$arr[57,11] = 6.0
$arr[57,12] = 1.333
$arr[57,13] = 5.0
$arr[57,14] = 6.0
$arr[57,15] = -1.5
$arr[57,16] = 1.85
$arr[57,17] = 1.95
$arr[57,18] = 3.5
$arr[57,19] = 2.0
$arr[57,20] = 1.8
$arr[57,21] = 0.333
$arr[57,22] = -1.0
$arr[57,23] = -1.0
$arr[57,24] = 0.8500000000000001
$arr[57,25] = -1.0
$arr[57,26] = -1.0
$arr[57,27] = 0.8
$arr[58,0] = 58
$arr[58,1] = 7873625
$arr[58,2] = 'Germany Landesliga'
$arr[58,3] = 45346.64583333334
$arr[58,4] = 'SG 2000 MulheimKarlich'
$arr[58,5] = 'FSG EhrangPfalzel'
$arr[58,6] = 1
$arr[58,7] = 0
$arr[58,8] = 'H'
$arr[58,9] = 2.0
$arr[58,10] = 3.6
$arr[58,11] = 3.0
$arr[58,12] = 1.95
$arr[58,13] = 3.75
$arr[58,14] = 3.0
$arr[58,15] = -0.25
$arr[58,16] = 1.8
$arr[58,17] = 2.0
$arr[58,18] = 3.25
$arr[58,19] = 1.825
$arr[58,20] = 1.975
$arr[58,21] = 0.95
$arr[58,22] = -1.0
$arr[58,23] = -1.0
$arr[58,24] = 0.8
$arr[58,25] = -1.0
$arr[58,26] = -1.0
$arr[58,27] = 0.9750000000000001
$arr[59,0] = 59
$arr[59,1] = 7890341
$arr[59,2] = 'Germany Landesliga'
$arr[59,3] = 45350.625
$arr[59,4] = 'BSC Rapid Chemnitz'
$arr[59,5] = 'VfB Fortuna Chemnitz'
$arr[59,6] = 2
$arr[59,7] = 0
$arr[59,8] = 'H'
$arr[59,9] = 13.0
$arr[59,10] = 7.5
$arr[59,11] = 1.125
$arr[59,12] = 21.0
$arr[59,13] = 10.0
$arr[59,14] = 1.083
$arr[59,15] = 2.75
$arr[59,16] = 2.075
$arr[59,17] = 1.725
$arr[59,18] = 4.0
$arr[59,19] = 1.875
$arr[59,20] = 1.925
$arr[59,21] = 20.0
$arr[59,22] = -1.0
$arr[59,23] = -1.0
$arr[59,24] = 1.075
$arr[59,25] = -1.0
$arr[59,26] = -1.0
$arr[59,27] = 0.925
$arr[60,0] = 60
$arr[60,1] = 7897140
$arr[60,2] = 'Germany Landesliga'
$arr[60,3] = 45352.66666666666
$arr[60,4] = 'SG Unterrath'
$arr[60,5] = 'ASV Suchteln'
$arr[60,6] = 0
$arr[60,7] = 3
$arr[60,8] = 'A'
$arr[60,9] = 2.15
$arr[60,10] = 4.2
$arr[60,11] = 2.4
$arr[60,12] = 2.15
$arr[60,13] = 4.2
$arr[60,14] = 2.4
$arr[60,15] = 0.0
$arr[60,16] = 1.8
$arr[60,17] = 2.0
$arr[60,18] = 3.75
$arr[60,19] = 1.85
$arr[60,20] = 1.95
$arr[60,21] = -1.0
$arr[60,22] = -1.0
$arr[60,23] = 1.4
$arr[60,24] = -1.0
$arr[60,25] = 1.0
$arr[60,26] = -1.0
$arr[60,27] = 0.95
$arr[61,0] = 61
$arr[61,1] = 7905567
$arr[61,2] = 'Germany Landesliga'
$arr[61,3] = 45354.46875
$arr[61,4] = 'SV Schlebusch'
$arr[61,5] = 'SC Rheinbach'
$arr[61,6] = 0
$arr[61,7] = 0
$arr[61,8] = 'D'
$arr[61,9] = 2.9
$arr[61,10] = 3.4
$arr[61,11] = 2.1
$arr[61,12] = 2.9
$arr[61,13] = 3.4
$arr[61,14] = 2.1
$arr[61,15] = 0.25
$arr[61,16] = 1.9
$arr[61,17] = 1.9
$arr[61,18] = 3.0
$arr[61,19] = 1.85
$arr[61,20] = 1.95
$arr[61,21] = -1.0
$arr[61,22] = 2.4
$arr[61,23] = -1.0
$arr[61,24] = 0.45
$arr[61,25] = -0.5
$arr[61,26] = -1.0
$arr[61,27] = 0.95
$arr[62,0] = 62
$arr[62,1] = 7930369
$arr[62,2] = 'Germany Landesliga'
$arr[62,3] = 45359.64583333334
$arr[62,4] = 'MSV Dusseldorf'
$arr[62,5] = 'SC Victoria Mennrath'
$arr[62,6] = 0
$arr[62,7] = 4
$arr[62,8] = 'A'
$arr[62,9] = 2.3
$arr[62,10] = 3.6
$arr[62,11] = 2.5
$arr[62,12] = 2.2
$arr[62,13] = 3.5
$arr[62,14] = 2.625
$arr[62,15] = -0.25
$arr[62,16] = 2.0
$arr[62,17] = 1.8
$arr[62,18] = 3.5
$arr[62,19] = 2.0
$arr[62,20] = 1.8
$arr[62,21] = -1.0
$arr[62,22] = -1.0
$arr[62,23] = 1.625
$arr[62,24] = -1.0
$arr[62,25] = 0.8
$arr[62,26] = 1.0
$arr[62,27] = -1.0
$arr[63,0] = 63
$arr[63,1] = 7940440
$arr[63,2] = 'Germany Landesliga'
$arr[63,3] = 45361.41666666666
$arr[63,4] = 'SSV Markranstadt'
$arr[63,5] = 'SG Taucha 99'
$arr[63,6] = 2
$arr[63,7] = 0
$arr[63,8] = 'H'
$arr[63,9] = 2.45
$arr[63,10] = 3.25
$arr[63,11] = 2.5
$arr[63,12] = 2.45
$arr[63,13] = 3.3
$arr[63,14] = 2.5
$arr[63,15] = 0.0
$arr[63,16] = 1.875
$arr[63,17] = 1.925
$arr[63,18] = 2.25
$arr[63,19] = 1.9
$arr[63,20] = 1.9
$arr[63,21] = 1.45
$arr[63,22] = -1.0
$arr[63,23] = -1.0
$arr[63,24] = 0.875
$arr[63,25] = -1.0
$arr[63,26] = -0.5
$arr[63,27] = 0.45
$arr[64,0] = 64
$arr[64,1] = 7951581
$arr[64,2] = 'Germany Landesliga'
$arr[64,3] = 45363.66666666666
$arr[64,4] = 'FC Concordia 03'
$arr[64,5] = 'Germania Schneiche'
$arr[64,6] = 0
$arr[64,7] = 3
$arr[64,8] = 'A'
$arr[64,9] = 2.9
$arr[64,10] = 4.0
$arr[64,11] = 1.909
$arr[64,12] = 2.875
$arr[64,13] = 4.0
$arr[64,14] = 1.909
$arr[64,15] = 0.5
$arr[64,16] = 1.825
$arr[64,17] = 1.975
$arr[64,18] = 3.75
$arr[64,19] = 1.925
$arr[64,20] = 1.875
$arr[64,21] = -1.0
$arr[64,22] = -1.0
$arr[64,23] = 0.909
$arr[64,24] = -1.0
$arr[64,25] = 0.9750000000000001
$arr[64,26] = -1.0
$arr[64,27] = 0.875
$arr[65,0] = 65
$arr[65,1] = 7950775
$arr[65,2] = 'Germany Landesliga'
$arr[65,3] = 45364.66666666666
$arr[65,4] = 'TV Dinklage'
$arr[65,5] = 'Viktoria Gesmold'
$arr[65,6] = 4
$arr[65,7] = 0
$arr[65,8] = 'H'
$arr[65,9] = 2.15
$arr[65,10] = 3.4
$arr[65,11] = 2.8
$arr[65,12] = 2.25
$arr[65,13] = 3.5
$arr[65,14] = 2.625
$arr[65,15] = 0.0
$arr[65,16] = 1.725
$arr[65,17] = 2.075
$arr[65,18] = 3.0
$arr[65,19] = 1.8
$arr[65,20] = 2.0
$arr[65,21] = 1.25
$arr[65,22] = -1.0
$arr[65,23] = -1.0
$arr[65,24] = 0.7250000000000001
$arr[65,25] = -1.0
$arr[65,26] = 0.8
$arr[65,27] = -1.0
$arr[66,0] = 66
$arr[66,1] = 7969732
$arr[66,2] = 'Germany Landesliga'
$arr[66,3] = 45368.45833333334
$arr[66,4] = 'FC Monheim'
$arr[66,5] = 'TuRU Dsseldorf'
$arr[66,6] = 2
$arr[66,7] = 1
$arr[66,8] = 'H'
$arr[66,9] = 1.333
$arr[66,10] = 5.0
$arr[66,11] = 6.0
$arr[66,12] = 1.285
$arr[66,13] = 5.0
$arr[66,14] = 7.5
$arr[66,15] = -1.75
$arr[66,16] = 1.975
$arr[66,17] = 1.825
$arr[66,18] = 3.25
$arr[66,19] = 2.025
$arr[66,20] = 1.775
$arr[66,21] = 0.2849999999999999
$arr[66,22] = -1.0
$arr[66,23] = -1.0
$arr[66,24] = -1.0
$arr[66,25] = 0.825
$arr[66,26] = -0.5
$arr[66,27] = 0.3875
$arr[67,0] = 67
$arr[67,1] = 7969733
$arr[67,2] = 'Germany Landesliga'
$arr[67,3] = 45368.47916666666
$arr[67,4] = 'SC Dsseldorf West'
$arr[67,5] = 'DV Solingen'
$arr[67,6] = 0
$arr[67,7] = 2
$arr[67,8] = 'A'
$arr[67,9] = 3.9
$arr[67,10] = 4.5
$arr[67,11] = 1.571
$arr[67,12] = 3.1
$arr[67,13] = 4.2
$arr[67,14] = 1.8
$arr[67,15] = 0.5
$arr[67,16] = 1.975
$arr[67,17] = 1.825
$arr[67,18] = 3.25
$arr[67,19] = 1.975
$arr[67,20] = 1.825
$arr[67,21] = -1.0
$arr[67,22] = -1.0
$arr[67,23] = 0.8
$arr[67,24] = -1.0
$arr[67,25] = 0.825
$arr[67,26] = -1.0
$arr[67,27] = 0.825
$arr[68,0] = 68
$arr[68,1] = 7990732
$arr[68,2] = 'Germany Landesliga'
$arr[68,3] = 45373.64583333334
$arr[68,4] = 'TuRU Dsseldorf'
$arr[68,5] = 'FC Viersen'
$arr[68,6] = 1
$arr[68,7] = 1
$arr[68,8] = 'D'
$arr[68,9] = 1.909
$arr[68,10] = 3.25
$arr[68,11] = 3.5
$arr[68,12] = 2.05
$arr[68,13] = 3.2
$arr[68,14] = 3.3
$arr[68,15] = -0.25
$arr[68,16] = 1.8
$arr[68,17] = 2.0
$arr[68,18] = 2.5
$arr[68,19] = 1.8
$arr[68,20] = 2.0
$arr[68,21] = -1.0
$arr[68,22] = 2.2
$arr[68,23] = -1.0
$arr[68,24] = -0.5
$arr[68,25] = 0.5
$arr[68,26] = -1.0
$arr[68,27] = 1.0
$arr[69,0] = 69
$arr[69,1] = 7992020
$arr[69,2] = 'Germany Landesliga'
$arr[69,3] = 45373.66666666666
$arr[69,4] = 'ASV Suchteln'
$arr[69,5] = 'SC Dsseldorf West'
$arr[69,6] = 4
$arr[69,7] = 0
$arr[69,8] = 'H'
$arr[69,9] = 1.285
$arr[69,10] = 5.5
$arr[69,11] = 6.5
$arr[69,12] = 1.444
$arr[69,13] = 4.75
$arr[69,14] = 4.75
$arr[69,15] = -1.25
$arr[69,16] = 1.875
$arr[69,17] = 1.925
$arr[69,18] = 3.5
$arr[69,19] = 1.95
$arr[69,20] = 1.85
$arr[69,21] = 0.444
$arr[69,22] = -1.0
$arr[69,23] = -1.0
$arr[69,24] = 0.875
$arr[69,25] = -1.0
$arr[69,26] = 0.95
$arr[69,27] = -1.0
$arr[70,0] = 70
$arr[70,1] = 7999992
$arr[70,2] = 'Germany Landesliga'
$arr[70,3] = 45375.47916666666
$arr[70,4] = 'MSV Dusseldorf'
$arr[70,5] = 'FC Remscheid'
$arr[70,6] = 3
$arr[70,7] = 3
$arr[70,8] = 'D'
$arr[70,9] = 4.5
$arr[70,10] = 4.5
$arr[70,11] = 1.5
$arr[70,12] = 4.2
$arr[70,13] = 4.2
$arr[70,14] = 1.6
$arr[70,15] = 1.0
$arr[70,16] = 1.8
$arr[70,17] = 2.0
$arr[70,18] = 3.25
$arr[70,19] = 2.0
$arr[70,20] = 1.8
$arr[70,21] = -1.0
$arr[70,22] = 3.2
$arr[70,23] = -1.0
$arr[70,24] = 0.8
$arr[70,25] = -1.0
$arr[70,26] = 1.0
$arr[70,27] = -1.0
$arr[71,0] = 71
$arr[71,1] = 8013545
$arr[71,2] = 'Germany Landesliga'
$arr[71,3] = 45378.64583333334
$arr[71,4] = 'TuRU Dsseldorf'
$arr[71,5] = 'MSV Dusseldorf'
$arr[71,6] = 2
$arr[71,7] = 0
$arr[71,8] = 'H'
$arr[71,9] = 1.333
$arr[71,10] = 5.0
$arr[71,11] = 6.0
$arr[71,12] = 1.363
$arr[71,13] = 4.75
$arr[71,14] = 6.0
$arr[71,15] = -1.5
$arr[71,16] = 1.9
$arr[71,17] = 1.9
$arr[71,18] = 3.25
$arr[71,19] = 1.8
$arr[71,20] = 2.0
$arr[71,21] = 0.363
$arr[71,22] = -1.0
$arr[71,23] = -1.0
$arr[71,24] = 0.8999999999999999
$arr[71,25] = -1.0
$arr[71,26] = -1.0
$arr[71,27] = 1.0
$arr[72,0] = 72
$arr[72,1] = 8013546
$arr[72,2] = 'Germany Landesliga'
$arr[72,3] = 45378.65625
$arr[72,4] = 'ASV Mettmann'
$arr[72,5] = 'DJK Dilkrath'
$arr[72,6] = 0
$arr[72,7] = 1
$arr[72,8] = 'A'
$arr[72,9] = 2.75
$arr[72,10] = 4.0
$arr[72,11] = 2.0
$arr[72,12] = 3.0
$arr[72,13] = 4.0
$arr[72,14] = 1.909
$arr[72,15] = 0.5
$arr[72,16] = 1.85
$arr[72,17] = 1.95
$arr[72,18] = 4.0
$arr[72,19] = 1.925
$arr[72,20] = 1.875
$arr[72,21] = -1.0
$arr[72,22] = -1.0
$arr[72,23] = 0.909
$arr[72,24] = -1.0
$arr[72,25] = 0.95
$arr[72,26] = -1.0
$arr[72,27] = 0.875
$arr[73,0] = 73
$arr[73,1] = 8017334
$arr[73,2] = 'Germany Landesliga'
$arr[73,3] = 45379.625
$arr[73,4] = 'VfB Fortuna Chemnitz'
$arr[73,5] = 'SG Handwerk Rabenstein'
$arr[73,6] = 0
$arr[73,7] = 2
$arr[73,8] = 'A'
$arr[73,9] = 4.0
$arr[73,10] = 4.0
$arr[73,11] = 1.615
$arr[73,12] = 4.0
$arr[73,13] = 4.0
$arr[73,14] = 1.615
$arr[73,15] = 0.75
$arr[73,16] = 1.975
$arr[73,17] = 1.825
$arr[73,18] = 2.75
$arr[73,19] = 1.975
$arr[73,20] = 1.825
$arr[73,21] = -1.0
$arr[73,22] = -1.0
$arr[73,23] = 0.615
$arr[73,24] = -1.0
$arr[73,25] = 0.825
$arr[73,26] = -1.0
$arr[73,27] = 0.825
$arr[74,0] = 74
$arr[74,1] = 8034132
$arr[74,2] = 'Germany Landesliga'
$arr[74,3] = 45383.41666666666
$arr[74,4] = 'SC Fortuna Bonn'
$arr[74,5] = 'FSV NeunkirchenSeelscheid'
$arr[74,6] = 2
$arr[74,7] = 2
$arr[74,8] = 'D'
$arr[74,9] = 2.2
$arr[74,10] = 3.6
$arr[74,11] = 2.625
$arr[74,12] = 2.25
$arr[74,13] = 3.6
$arr[74,14] = 2.55
$arr[74,15] = 0.0
$arr[74,16] = 1.775
$arr[74,17] = 2.025
$arr[74,18] = 3.0
$arr[74,19] = 1.8
$arr[74,20] = 2.0
$arr[74,21] = -1.0
$arr[74,22] = 2.6
$arr[74,23] = -1.0
$arr[74,24] = 0.0
$arr[74,25] = 0.0
$arr[74,26] = 0.8
$arr[74,27] = -1.0
$arr[75,0] = 75
$arr[75,1] = 8039588
$arr[75,2] = 'Germany Landesliga'
$arr[75,3] = 45383.60416666666
$arr[75,4] = 'SV Vorwarts Nordhorn'
$arr[75,5] = 'TV Dinklage'
$arr[75,6] = 2
$arr[75,7] = 2
$arr[75,8] = 'D'
$arr[75,9] = 2.0
$arr[75,10] = 3.75
$arr[75,11] = 2.875
$arr[75,12] = 2.0
$arr[75,13] = 3.75
$arr[75,14] = 2.9
$arr[75,15] = -0.25
$arr[75,16] = 1.8
$arr[75,17] = 2.0
$arr[75,18] = 3.25
$arr[75,19] = 1.875
$arr[75,20] = 1.925
$arr[75,21] = -1.0
$arr[75,22] = 2.75
$arr[75,23] = -1.0
$arr[75,24] = -0.5
$arr[75,25] = 0.5
$arr[75,26] = 0.875
$arr[75,27] = -1.0
$arr[76,0] = 76
$arr[76,1] = 8050870
$arr[76,2] = 'Germany Landesliga'
$arr[76,3] = 45387.58333333334
$arr[76,4] = 'FC Kray'
$arr[76,5] = 'Spvgg Steele 0309'
$arr[76,6] = 0
$arr[76,7] = 0
$arr[76,8] = 'D'
$arr[76,9] = 1.8
$arr[76,10] = 3.8
$arr[76,11] = 3.4
$arr[76,12] = 1.8
$arr[76,13] = 3.75
$arr[76,14] = 3.5
$arr[76,15] = -0.5
$arr[76,16] = 1.825
$arr[76,17] = 1.975
$arr[76,18] = 3.25
$arr[76,19] = 1.875
$arr[76,20] = 1.925
$arr[76,21] = -1.0
$arr[76,22] = 2.75
$arr[76,23] = -1.0
$arr[76,24] = -1.0
$arr[76,25] = 0.9750000000000001
$arr[76,26] = -1.0
$arr[76,27] = 0.925
$arr[77,0] = 77
$arr[77,1] = 8050971
$arr[77,2] = 'Germany Landesliga'
$arr[77,3] = 45387.61458333334
$arr[77,4] = 'TuS Kirchberg 1909'
$arr[77,5] = 'FSV Salmrohr'
$arr[77,6] = 0
$arr[77,7] = 0
$arr[77,8] = 'D'
$arr[77,9] = 1.909
$arr[77,10] = 3.75
$arr[77,11] = 3.1
$arr[77,12] = 1.909
$arr[77,13] = 3.75
$arr[77,14] = 3.2
$arr[77,15] = -0.5
$arr[77,16] = 1.95
$arr[77,17] = 1.85
$arr[77,18] = 4.0
$arr[77,19] = 2.0
$arr[77,20] = 1.8
$arr[77,21] = -1.0
$arr[77,22] = 2.75
$arr[77,23] = -1.0
$arr[77,24] = -1.0
$arr[77,25] = 0.8500000000000001
$arr[77,26] = -1.0
$arr[77,27] = 0.8
$arr[78,0] = 78
$arr[78,1] = 8059643
$arr[78,2] = 'Germany Landesliga'
$arr[78,3] = 45389.4375
$arr[78,4] = 'SC Victoria Mennrath'
$arr[78,5] = 'SG Unterrath'
$arr[78,6] = 0
$arr[78,7] = 0
$arr[78,8] = 'D'
$arr[78,9] = 2.2
$arr[78,10] = 3.5
$arr[78,11] = 2.7
$arr[78,12] = 2.2
$arr[78,13] = 3.6
$arr[78,14] = 2.625
$arr[78,15] = -0.25
$arr[78,16] = 2.0
$arr[78,17] = 1.8
$arr[78,18] = 3.25
$arr[78,19] = 1.8
$arr[78,20] = 2.0
$arr[78,21] = -1.0
$arr[78,22] = 2.6
$arr[78,23] = -1.0
$arr[78,24] = -0.5
$arr[78,25] = 0.4
$arr[78,26] = -1.0
$arr[78,27] = 1.0
$arr[79,0] = 79
$arr[79,1] = 8059644
$arr[79,2] = 'Germany Landesliga'
$arr[79,3] = 45389.4375
$arr[79,4] = 'ASV Suchteln'
$arr[79,5] = 'FC Monheim'
$arr[79,6] = 2
$arr[79,7] = 4
$arr[79,8] = 'A'
$arr[79,9] = 3.75
$arr[79,10] = 4.0
$arr[79,11] = 1.666
$arr[79,12] = 4.75
$arr[79,13] = 4.2
$arr[79,14] = 1.5
$arr[79,15] = 1.0
$arr[79,16] = 1.975
$arr[79,17] = 1.825
$arr[79,18] = 3.5
$arr[79,19] = 1.9
$arr[79,20] = 1.9
$arr[79,21] = -1.0
$arr[79,22] = -1.0
$arr[79,23] = 0.5
$arr[79,24] = -1.0
$arr[79,25] = 0.825
$arr[79,26] = 0.8999999999999999
$arr[79,27] = -1.0
$arr[80,0] = 80
$arr[80,1] = 8075592
$arr[80,2] = 'Germany Landesliga'
$arr[80,3] = 45392.58333333334
$arr[80,4] = 'SGS EssenSchonebeck'
$arr[80,5] = 'DJK Arminia Klosterhardt'
$arr[80,6] = 1
$arr[80,7] = 1
$arr[80,8] = 'D'
$arr[80,9] = 2.4
$arr[80,10] = 3.4
$arr[80,11] = 2.5
$arr[80,12] = 2.4
$arr[80,13] = 3.4
$arr[80,14] = 2.5
$arr[80,15] = 0.0
$arr[80,16] = 1.85
$arr[80,17] = 1.95
$arr[80,18] = 3.25
$arr[80,19] = 2.0
$arr[80,20] = 1.8
$arr[80,21] = -1.0
$arr[80,22] = 2.4
$arr[80,23] = -1.0
$arr[80,24] = 0.0
$arr[80,25] = 0.0
$arr[80,26] = -1.0
$arr[80,27] = 0.8
$arr[81,0] = 81
$arr[81,1] = 8075296
$arr[81,2] = 'Germany Landesliga'
$arr[81,3] = 45392.60416666666
$arr[81,4] = 'FC Monheim'
$arr[81,5] = 'VFB Hilden II'
$arr[81,6] = 1
$arr[81,7] = 2
$arr[81,8] = 'A'
$arr[81,9] = 1.533
$arr[81,10] = 4.75
$arr[81,11] = 4.0
$arr[81,12] = 1.4
$arr[81,13] = 5.25
$arr[81,14] = 5.0
$arr[81,15] = -1.5
$arr[81,16] = 1.975
$arr[81,17] = 1.825
$arr[81,18] = 3.75
$arr[81,19] = 1.9
$arr[81,20] = 1.9
$arr[81,21] = -1.0
$arr[81,22] = -1.0
$arr[81,23] = 4.0
$arr[81,24] = -1.0
$arr[81,25] = 0.825
$arr[81,26] = -1.0
$arr[81,27] = 0.8999999999999999
$arr[82,0] = 82
$arr[82,1] = 8075530
$arr[82,2] = 'Germany Landesliga'
$arr[82,3] = 45392.60416666666
$arr[82,4] = 'TuRU Dsseldorf'
$arr[82,5] = 'DV Solingen'
$arr[82,6] = 2
$arr[82,7] = 0
$arr[82,8] = 'H'
$arr[82,9] = 2.1
$arr[82,10] = 3.75
$arr[82,11] = 2.7
$arr[82,12] = 2.375
$arr[82,13] = 3.75
$arr[82,14] = 2.45
$arr[82,15] = 0.0
$arr[82,16] = 1.85
$arr[82,17] = 1.95
$arr[82,18] = 3.0
$arr[82,19] = 1.85
$arr[82,20] = 1.95
$arr[82,21] = 1.375
$arr[82,22] = -1.0
$arr[82,23] = -1.0
$arr[82,24] = 0.8500000000000001
$arr[82,25] = -1.0
$arr[82,26] = -1.0
$arr[82,27] = 0.95
$arr[83,0] = 83
$arr[83,1] = 8075670
$arr[83,2] = 'Germany Landesliga'
$arr[83,3] = 45392.625
$arr[83,4] = 'SF Niederwenigern'
$arr[83,5] = 'VfB Frohnhausen'
$arr[83,6] = 2
$arr[83,7] = 1
$arr[83,8] = 'H'
$arr[83,9] = 1.222
$arr[83,10] = 6.5
$arr[83,11] = 9.5
$arr[83,12] = 1.181
$arr[83,13] = 7.0
$arr[83,14] = 11.0
$arr[83,15] = -2.5
$arr[83,16] = 1.925
$arr[83,17] = 1.875
$arr[83,18] = 4.25
$arr[83,19] = 1.8
$arr[83,20] = 2.0
$arr[83,21] = 0.181
$arr[83,22] = -1.0
$arr[83,23] = -1.0
$arr[83,24] = -1.0
$arr[83,25] = 0.875
$arr[83,26] = -1.0
$arr[83,27] = 1.0
$arr[84,0] = 84
$arr[84,1] = 8075593
$arr[84,2] = 'Germany Landesliga'
$arr[84,3] = 45392.625
$arr[84,4] = 'FSV Duisburg'
$arr[84,5] = 'ESC Rellinghausen'
$arr[84,6] = 0
$arr[84,7] = 2
$arr[84,8] = 'A'
$arr[84,9] = 6.5
$arr[84,10] = 4.5
$arr[84,11] = 1.363
$arr[84,12] = 4.75
$arr[84,13] = 4.5
$arr[84,14] = 1.5
$arr[84,15] = 1.25
$arr[84,16] = 1.875
$arr[84,17] = 1.925
$arr[84,18] = 3.75
$arr[84,19] = 2.0
$arr[84,20] = 1.8
$arr[84,21] = -1.0
$arr[84,22] = -1.0
$arr[84,23] = 0.5
$arr[84,24] = -1.0
$arr[84,25] = 0.925
$arr[84,26] = -1.0
$arr[84,27] = 0.8
$arr[85,0] = 85
$arr[85,1] = 8076436
$arr[85,2] = 'Germany Landesliga'
$arr[85,3] = 45393.58333333334
$arr[85,4] = 'Spvgg Steele 0309'
$arr[85,5] = 'VfB Speldorf'
$arr[85,6] = 1
$arr[85,7] = 3
$arr[85,8] = 'A'
$arr[85,9] = 2.7
$arr[85,10] = 3.75
$arr[85,11] = 2.1
$arr[85,12] = 2.75
$arr[85,13] = 3.75
$arr[85,14] = 2.1
$arr[85,15] = 0.25
$arr[85,16] = 1.875
$arr[85,17] = 1.925
$arr[85,18] = 3.75
$arr[85,19] = 1.975
$arr[85,20] = 1.825
$arr[85,21] = -1.0
$arr[85,22] = -1.0
$arr[85,23] = 1.1
$arr[85,24] = -1.0
$arr[85,25] = 0.925
$arr[85,26] = 0.4875
$arr[85,27] = -0.5
$arr[86,0] = 86
$arr[86,1] = 8076438
$arr[86,2] = 'Germany Landesliga'
$arr[86,3] = 45393.60416666666
$arr[86,4] = 'SV Rott 1927'
$arr[86,5] = 'SV Breinig'
$arr[86,6] = 1
$arr[86,7] = 4
$arr[86,8] = 'A'
$arr[86,9] = 1.363
$arr[86,10] = 5.0
$arr[86,11] = 5.5
$arr[86,12] = 1.65
$arr[86,13] = 4.5
$arr[86,14] = 3.5
$arr[86,15] = -0.75
$arr[86,16] = 1.825
$arr[86,17] = 1.975
$arr[86,18] = 3.5
$arr[86,19] = 1.975
$arr[86,20] = 1.825
$arr[86,21] = -1.0
$arr[86,22] = -1.0
$arr[86,23] = 2.5
$arr[86,24] = -1.0
$arr[86,25] = 0.9750000000000001
$arr[86,26] = 0.9750000000000001
$arr[86,27] = -1.0
$arr[87,0] = 87
$arr[87,1] = 8076477
$arr[87,2] = 'Germany Landesliga'
$arr[87,3] = 45393.60416666666
$arr[87,4] = 'BSC Hastedt'
$arr[87,5] = 'SC Weyhe'
$arr[87,6] = 2
$arr[87,7] = 2
$arr[87,8] = 'D'
$arr[87,9] = 1.833
$arr[87,10] = 4.0
$arr[87,11] = 3.1
$arr[87,12] = 1.833
$arr[87,13] = 4.0
$arr[87,14] = 3.1
$arr[87,15] = -0.5
$arr[87,16] = 1.875
$arr[87,17] = 1.925
$arr[87,18] = 3.75
$arr[87,19] = 1.975
$arr[87,20] = 1.825
$arr[87,21] = -1.0
$arr[87,22] = 3.0
$arr[87,23] = -1.0
$arr[87,24] = -1.0
$arr[87,25] = 0.925
$arr[87,26] = 0.4875
$arr[87,27] = -0.5
$arr[88,0] = 88
$arr[88,1] = 8077795
$arr[88,2] = 'Germany Landesliga'
$arr[88,3] = 45393.60416666666
$arr[88,4] = 'Eiche Horn'
$arr[88,5] = 'SVGO Bremen'
$arr[88,6] = 5
$arr[88,7] = 3
$arr[88,8] = 'H'
$arr[88,9] = 1.142
$arr[88,10] = 7.0
$arr[88,11] = 10.0
$arr[88,12] = 1.083
$arr[88,13] = 11.0
$arr[88,14] = 19.0
$arr[88,15] = -3.5
$arr[88,16] = 1.975
$arr[88,17] = 1.825
$arr[88,18] = 5.0
$arr[88,19] = 1.825
$arr[88,20] = 1.975
$arr[88,21] = 0.08299999999999996
$arr[88,22] = -1.0
$arr[88,23] = -1.0
$arr[88,24] = -1.0
$arr[88,25] = 0.825
$arr[88,26] = 0.825
$arr[88,27] = -1.0
$arr[89,0] = 89
$arr[89,1] = 8076435
$arr[89,2] = 'Germany Landesliga'
$arr[89,3] = 45393.625
$arr[89,4] = 'Cronenberger SC'
$arr[89,5] = 'SC Victoria Mennrath'
$arr[89,6] = 2
$arr[89,7] = 5
$arr[89,8] = 'A'
$arr[89,9] = 1.727
$arr[89,10] = 4.0
$arr[89,11] = 3.5
$arr[89,12] = 1.7
$arr[89,13] = 4.2
$arr[89,14] = 3.5
$arr[89,15] = -0.75
$arr[89,16] = 1.925
$arr[89,17] = 1.875
$arr[89,18] = 3.25
$arr[89,19] = 1.9
$arr[89,20] = 1.9
$arr[89,21] = -1.0
$arr[89,22] = -1.0
$arr[89,23] = 2.5
$arr[89,24] = -1.0
$arr[89,25] = 0.875
$arr[89,26] = 0.8999999999999999
$arr[89,27] = -1.0
$arr[90,0] = 90
$arr[90,1] = 8086837
$arr[90,2] = 'Germany Landesliga'
$arr[90,3] = 45394.60416666666
$arr[90,4] = 'TuS Immendorf'
$arr[90,5] = 'Sportfreunde Eisbachtal'
$arr[90,6] = 2
$arr[90,7] = 2
$arr[90,8] = 'D'
$arr[90,9] = 5.0
$arr[90,10] = 5.5
$arr[90,11] = 1.363
$arr[90,12] = 4.0
$arr[90,13] = 5.5
$arr[90,14] = 1.45
$arr[90,15] = 1.25
$arr[90,16] = 1.85
$arr[90,17] = 1.95
$arr[90,18] = 3.5
$arr[90,19] = 1.975
$arr[90,20] = 1.825
$arr[90,21] = -1.0
$arr[90,22] = 4.5
$arr[90,23] = -1.0
$arr[90,24] = 0.8500000000000001
$arr[90,25] = -1.0
$arr[90,26] = 0.9750000000000001
$arr[90,27] = -1.0
$arr[91,0] = 91
$arr[91,1] = 8121115
$arr[91,2] = 'Germany Landesliga'
$arr[91,3] = 45403.42708333334
$arr[91,4] = 'SV Schlebusch'
$arr[91,5] = 'BlauWei Friesdorf'
$arr[91,6] = 4
$arr[91,7] = 2
$arr[91,8] = 'H'
$arr[91,9] = 1.833
$arr[91,10] = 3.5
$arr[91,11] = 3.5
$arr[91,12] = 1.8
$arr[91,13] = 3.6
$arr[91,14] = 3.5
$arr[91,15] = -0.5
$arr[91,16] = 1.875
$arr[91,17] = 1.925
$arr[91,18] = 3.25
$arr[91,19] = 1.85
$arr[91,20] = 1.95
$arr[91,21] = 0.8
$arr[91,22] = -1.0
$arr[91,23] = -1.0
$arr[91,24] = 0.875
$arr[91,25] = -1.0
$arr[91,26] = 0.8500000000000001
$arr[91,27] = -1.0
$arr[92,0] = 92
$arr[92,1] = 8140569
$arr[92,2] = 'Germany Landesliga'
$arr[92,3] = 45407.625
$arr[92,4] = 'MTV Soltau'
$arr[92,5] = 'TSV Elstorf'
$arr[92,6] = 5
$arr[92,7] = 3
$arr[92,8] = 'H'
$arr[92,9] = 1.4
$arr[92,10] = 4.5
$arr[92,11] = 5.5
$arr[92,12] = 1.363
$arr[92,13] = 4.5
$arr[92,14] = 5.75
$arr[92,15] = -1.5
$arr[92,16] = 1.875
$arr[92,17] = 1.925
$arr[92,18] = 4.0
$arr[92,19] = 1.9
$arr[92,20] = 1.9
$arr[92,21] = 0.363
$arr[92,22] = -1.0
$arr[92,23] = -1.0
$arr[92,24] = 0.875
$arr[92,25] = -1.0
$arr[92,26] = 0.8999999999999999
$arr[92,27] = -1.0
$arr[93,0] = 93
$arr[93,1] = 8143641
$arr[93,2] = 'Germany Landesliga'
$arr[93,3] = 45408.625
$arr[93,4] = 'SpVg Cologne Flittard'
$arr[93,5] = 'SSV Homburg Numbrecht'
$arr[93,6] = 0
$arr[93,7] = 0
$arr[93,8] = 'D'
$arr[93,9] = 2.2
$arr[93,10] = 4.2
$arr[93,11] = 2.375
$arr[93,12] = 2.2
$arr[93,13] = 4.2
$arr[93,14] = 2.375
$arr[93,15] = 0.0
$arr[93,16] = 1.825
$arr[93,17] = 1.975
$arr[93,18] = 3.5
$arr[93,19] = 1.825
$arr[93,20] = 1.975
$arr[93,21] = -1.0
$arr[93,22] = 3.2
$arr[93,23] = -1.0
$arr[93,24] = 0.0
$arr[93,25] = 0.0
$arr[93,26] = -1.0
$arr[93,27] = 0.9750000000000001
$arr[94,0] = 94
$arr[94,1] = 8163122
$arr[94,2] = 'Germany Landesliga'
$arr[94,3] = 45412.61111111111
$arr[94,4] = 'SG Union Klosterfelde'
$arr[94,5] = 'TSG Einheit Bernau'
$arr[94,6] = 2
$arr[94,7] = 0
$arr[94,8] = 'H'
$arr[94,9] = 1.6
$arr[94,10] = 4.2
$arr[94,11] = 4.0
$arr[94,12] = 1.6
$arr[94,13] = 4.2
$arr[94,14] = 4.0
$arr[94,15] = -1.0
$arr[94,16] = 2.025
$arr[94,17] = 1.775
$arr[94,18] = 3.0
$arr[94,19] = 1.95
$arr[94,20] = 1.85
$arr[94,21] = 0.6000000000000001
$arr[94,22] = -1.0
$arr[94,23] = -1.0
$arr[94,24] = 1.025
$arr[94,25] = -1.0
$arr[94,26] = -1.0
$arr[94,27] = 0.8500000000000001
$arr[95,0] = 95
$arr[95,1] = 8158915
$arr[95,2] = 'Germany Landesliga'
$arr[95,3] = 45412.625
$arr[95,4] = 'SV Budberg'
$arr[95,5] = 'FC Kray'
$arr[95,6] = 4
$arr[95,7] = 2
$arr[95,8] = 'H'
$arr[95,9] = 1.8
$arr[95,10] = 3.6
$arr[95,11] = 3.6
$arr[95,12] = 1.8
$arr[95,13] = 3.6
$arr[95,14] = 3.5
$arr[95,15] = -0.5
$arr[95,16] = 1.85
$arr[95,17] = 1.95
$arr[95,18] = 3.5
$arr[95,19] = 1.8
$arr[95,20] = 2.0
$arr[95,21] = 0.8
$arr[95,22] = -1.0
$arr[95,23] = -1.0
$arr[95,24] = 0.8500000000000001
$arr[95,25] = -1.0
$arr[95,26] = 0.8
$arr[95,27] = -1.0
$arr[96,0] = 96
$arr[96,1] = 8163123
$arr[96,2] = 'Germany Landesliga'
$arr[96,3] = 45412.625
$arr[96,4] = 'FC Concordia 03'
$arr[96,5] = 'SV Victoria Seelow'
$arr[96,6] = 5
$arr[96,7] = 1
$arr[96,8] = 'H'
$arr[96,9] = 1.7
$arr[96,10] = 3.75
$arr[96,11] = 3.75
$arr[96,12] = 1.7
$arr[96,13] = 3.75
$arr[96,14] = 3.75
$arr[96,15] = -0.5
$arr[96,16] = 1.75
$arr[96,17] = 2.05
$arr[96,18] = 3.75
$arr[96,19] = 1.85
$arr[96,20] = 1.95
$arr[96,21] = 0.7
$arr[96,22] = -1.0
$arr[96,23] = -1.0
$arr[96,24] = 0.75
$arr[96,25] = -1.0
$arr[96,26] = 0.8500000000000001
$arr[96,27] = -1.0
$arr[97,0] = 97
$arr[97,1] = 8161900
$arr[97,2] = 'Germany Landesliga'
$arr[97,3] = 45413.41666666666
$arr[97,4] = 'FC Monheim'
$arr[97,5] = 'SC Victoria Mennrath'
$arr[97,6] = 5
$arr[97,7] = 2
$arr[97,8] = 'H'
$arr[97,9] = 1.727
$arr[97,10] = 4.2
$arr[97,11] = 3.4
$arr[97,12] = 1.333
$arr[97,13] = 6.0
$arr[97,14] = 5.5
$arr[97,15] = -1.75
$arr[97,16] = 1.95
$arr[97,17] = 1.85
$arr[97,18] = 4.0
$arr[97,19] = 1.825
$arr[97,20] = 1.975
$arr[97,21] = 0.333
$arr[97,22] = -1.0
$arr[97,23] = -1.0
$arr[97,24] = 0.95
$arr[97,25] = -1.0
$arr[97,26] = 0.825
$arr[97,27] = -1.0
$arr[98,0] = 98
$arr[98,1] = 8166044
$arr[98,2] = 'Germany Landesliga'
$arr[98,3] = 45414.58333333334
$arr[98,4] = '1 FC Burg'
$arr[98,5] = 'SV Grohn'
$arr[98,6] = 1
$arr[98,7] = 2
$arr[98,8] = 'A'
$arr[98,9] = 2.25
$arr[98,10] = 3.5
$arr[98,11] = 2.6
$arr[98,12] = 2.25
$arr[98,13] = 4.0
$arr[98,14] = 2.4
$arr[98,15] = 0.0
$arr[98,16] = 1.825
$arr[98,17] = 1.975
$arr[98,18] = 5.0
$arr[98,19] = 1.925
$arr[98,20] = 1.875
$arr[98,21] = -1.0
$arr[98,22] = -1.0
$arr[98,23] = 1.4
$arr[98,24] = -1.0
$arr[98,25] = 0.9750000000000001
$arr[98,26] = -1.0
$arr[98,27] = 0.875
$arr[99,0] = 99
$arr[99,1] = 8166060
$arr[99,2] = 'Germany Landesliga'
$arr[99,3] = 45414.625
$arr[99,4] = 'SF Niederwenigern'
$arr[99,5] = 'SV Scherpenberg'
$arr[99,6] = 0
$arr[99,7] = 0
$arr[99,8] = 'D'
$arr[99,9] = 1.571
$arr[99,10] = 4.75
$arr[99,11] = 3.8
$arr[99,12] = 1.3
$arr[99,13] = 5.75
$arr[99,14] = 5.75
$arr[99,15] = -1.75
$arr[99,16] = 1.875
$arr[99,17] = 1.925
$arr[99,18] = 4.25
$arr[99,19] = 1.825
$arr[99,20] = 1.975
$arr[99,21] = -1.0
$arr[99,22] = 4.75
$arr[99,23] = -1.0
$arr[99,24] = -1.0
$arr[99,25] = 0.925
$arr[99,26] = -1.0
$arr[99,27] = 0.9750000000000001
$arr[100,0] = 100
$arr[100,1] = 8169516
$arr[100,2] = 'Germany Landesliga'
$arr[100,3] = 45415.58333333334
$arr[100,4] = 'SG Motor Wilsdruff'
$arr[100,5] = 'SSV Markranstadt'
$arr[100,6] = 0
$arr[100,7] = 0
$arr[100,8] = 'D'
$arr[100,9] = 3.0
$arr[100,10] = 3.6
$arr[100,11] = 2.0
$arr[100,12] = 3.3
$arr[100,13] = 3.75
$arr[100,14] = 1.833
$arr[100,15] = 0.5
$arr[100,16] = 1.9
$arr[100,17] = 1.9
$arr[100,18] = 3.0
$arr[100,19] = 1.875
$arr[100,20] = 1.925
$arr[100,21] = -1.0
$arr[100,22] = 2.75
$arr[100,23] = -1.0
$arr[100,24] = 0.8999999999999999
$arr[100,25] = -1.0
$arr[100,26] = -1.0
$arr[100,27] = 0.925
$arr[101,0] = 101
$arr[101,1] = 8169482
$arr[101,2] = 'Germany Landesliga'
$arr[101,3] = 45415.625
$arr[101,4] = 'SV Eintracht Verlautenheide'
$arr[101,5] = 'SC Germania Erftstadt'
$arr[101,6] = 3
$arr[101,7] = 5
$arr[101,8] = 'A'
$arr[101,9] = 1.571
$arr[101,10] = 4.0
$arr[101,11] = 4.333
$arr[101,12] = 1.4
$arr[101,13] = 4.5
$arr[101,14] = 5.5
$arr[101,15] = -1.5
$arr[101,16] = 1.925
$arr[101,17] = 1.875
$arr[101,18] = 4.0
$arr[101,19] = 1.825
$arr[101,20] = 1.975
$arr[101,21] = -1.0
$arr[101,22] = -1.0
$arr[101,23] = 4.5
$arr[101,24] = -1.0
$arr[101,25] = 0.875
$arr[101,26] = 0.825
$arr[101,27] = -1.0
$arr[102,0] = 102
$arr[102,1] = 8179934
$arr[102,2] = 'Germany Landesliga'
$arr[102,3] = 45417.41666666666
$arr[102,4] = 'Radefelder SV 90'
$arr[102,5] = 'Radebeuler BC 08'
$arr[102,6] = 2
$arr[102,7] = 2
$arr[102,8] = 'D'
$arr[102,9] = 2.3
$arr[102,10] = 3.2
$arr[102,11] = 2.75
$arr[102,12] = 2.5
$arr[102,13] = 3.2
$arr[102,14] = 2.5
$arr[102,15] = 0.0
$arr[102,16] = 1.9
$arr[102,17] = 1.9
$arr[102,18] = 2.5
$arr[102,19] = 1.9
$arr[102,20] = 1.9
$arr[102,21] = -1.0
$arr[102,22] = 2.2
$arr[102,23] = -1.0
$arr[102,24] = 0.0
$arr[102,25] = 0.0
$arr[102,26] = 0.8999999999999999
$arr[102,27] = -1.0
$arr[103,0] = 103
$arr[103,1] = 8179985
$arr[103,2] = 'Germany Landesliga'
$arr[103,3] = 45417.42708333334
$arr[103,4] = 'SV Schlebusch'
$arr[103,5] = 'SSV Homburg Numbrecht'
$arr[103,6] = 1
$arr[103,7] = 1
$arr[103,8] = 'D'
$arr[103,9] = 2.15
$arr[103,10] = 4.5
$arr[103,11] = 2.35
$arr[103,12] = 2.15
$arr[103,13] = 4.5
$arr[103,14] = 2.375
$arr[103,15] = -0.25
$arr[103,16] = 2.025
$arr[103,17] = 1.775
$arr[103,18] = 3.75
$arr[103,19] = 2.0
$arr[103,20] = 1.8
$arr[103,21] = -1.0
$arr[103,22] = 3.5
$arr[103,23] = -1.0
$arr[103,24] = -0.5
$arr[103,25] = 0.3875
$arr[103,26] = -1.0
$arr[103,27] = 0.8
$arr[104,0] = 104
$arr[104,1] = 8180064
$arr[104,2] = 'Germany Landesliga'
$arr[104,3] = 45417.4375
$arr[104,4] = 'Borussia LindenthalHohenlind'
$arr[104,5] = 'Kohlscheider BC'
$arr[104,6] = 6
$arr[104,7] = 1
$arr[104,8] = 'H'
$arr[104,9] = 1.3
$arr[104,10] = 5.25
$arr[104,11] = 6.5
$arr[104,12] = 1.3
$arr[104,13] = 5.5
$arr[104,14] = 6.0
$arr[104,15] = -1.75
$arr[104,16] = 1.875
$arr[104,17] = 1.925
$arr[104,18] = 4.0
$arr[104,19] = 1.85
$arr[104,20] = 1.95
$arr[104,21] = 0.3
$arr[104,22] = -1.0
$arr[104,23] = -1.0
$arr[104,24] = 0.875
$arr[104,25] = -1.0
$arr[104,26] = 0.8500000000000001
$arr[104,27] = -1.0
$ws.Range("A2").Resize(105,28).Value = $arr


# Copy styles for the newly added rows (96-106) from row 95 so that
# column A (bold/centered/bordered) and column D (date format) match
# the existing formatting pattern used throughout the sheet.
$ws.Range("A95").Copy() | Out-Null
$ws.Range("A96:A106").PasteSpecial(-4122) | Out-Null
$ws.Range("D95").Copy() | Out-Null
$ws.Range("D96:D106").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
